$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellB64 {
    param($addr, $b64text)
    $bytes = [System.Convert]::FromBase64String($b64text)
    $text = [System.Text.Encoding]::UTF8.GetString($bytes)
    $ws.Range($addr).Value = $text
}

# Remove the last data row (row 6) entirely - it is dropped in the new data set
$ws.Rows.Item(6).Delete()

# New columns added to the header row
Set-CellB64 "N1" "dWlfcmVwb3J0"
Set-CellB64 "O1" "YnJpZWY="
Set-CellB64 "P1" "cGFpbl9wb2ludF9yZXBvcnQ="
Set-CellB64 "Q1" "ZW1haWxfc2FtcGxl"

# ----- Row 2 -----
Set-CellB64 "A2" "Q2hJSlRjUms0UERJeG9rUjdlS2hrQzZ6U1Q4"
Set-CellB64 "B2" "TkogU21hbGwgQnVzaW5lc3MgRGV2ZWxvcG1lbnQgQ2VudGVyIEAgUnV0Z2Vycy1DYW1kZW4="
Set-CellB64 "C2" "Y29uc3VsdGFudCwgcG9pbnRfb2ZfaW50ZXJlc3QsIGVzdGFibGlzaG1lbnQ="
Set-CellB64 "D2" "KDg1NikgMjI1LTYyMjE="
Set-CellB64 "H2" "aHR0cDovL3JzYmRjLm9yZy8="
Set-CellB64 "I2" "aHR0cHM6Ly9tYXBzLmdvb2dsZS5jb20vP2NpZD00NTYwMzczMTEwMjY1NzMzODY5JmdfbXA9Q2lkbmIyOW5iR1V1YldGd2N5NXdiR0ZqWlhNdWRqRXVVR3hoWTJWekxsTmxZWEpqYUZSbGVIUVFBaGdFSUFB"
Set-CellB64 "K2" "QmVzdCBTbWFsbCBCdXNpbmVzcyByZXNvdXJjZSBhdmFpbGFibGUgZnJvbSBpZGVhIHRvIGV4cGFuc2lvbiBhbmQgYmV5b25kICEgfCBWZXJ5IGhlbHBmdWwgYW5kIGhpZ2hseSByZWNvbW1lbmRlZCBmb3Igc21hbGwgYnVzaW5lc3Mgb3duZXJzIG9yIHRob3NlIHdobyB3aXNoIHRvIHN0YXJ0IGEgYnVzaW5lc3MgfCBJbmZvIG9uIGhlbHAgeW91IHN0YXJ0IGEgYnVzaW5lc3MuIPCfkY0="
Set-CellB64 "L2" "cnNiZGNAY2FtZGVuLnJ1dGdlcnMuZWR1LCBhZzI5NTNAY2FtZGVuLnJ1dGdlcnMuZWR1"
Set-CellB64 "N2" "KipVSSBJbXByb3ZlbWVudHM6KioKMS4gKipCcmFuZGluZyBJbmNvbnNpc3RlbmN5OioqCiAgIC0gVGhlIFJ1dGdlcnMgYnJhbmRpbmcgaXMgc3Ryb25nIGluIHRoZSBoZWFkZXIsIGJ1dCB0aGUgdXNlIG9mIGRpZmZlcmVudCBjb2xvciBwYWxldHRlcyBmb3IgU0JEQywgbmF2aWdhdGlvbiwgYW5kIHRpdGxlcyBpbnRyb2R1Y2VzIHZpc3VhbCBmcmFnbWVudGF0aW9uLiBLZWVwaW5nIGEgdW5pZmllZCBjb2xvciBzY2hlbWUgd291bGQgY3JlYXRlIGEgbW9yZSBjb2hlc2l2ZSBicmFuZCBleHBlcmllbmNlLgoyLiAqKk91dGRhdGVkIFZpc3VhbCBTdHlsZToqKgogICAtIFRoZSBkZXNpZ24sIHdpdGggaXRzIGhlYXZ5IHVzZSBvZiBjb2xvcmVkIGJsb2NrcyBhbmQgZHJvcCBzaGFkb3dzLCBhbGlnbnMgbW9yZSB3aXRoIGVhcmx5IDIwMTBzIHRyZW5kcy4gTW9kZXJuIFVJcyBlbXBoYXNpemUgZmxhdCBkZXNpZ24sIHN1YnRsZSBncmFkaWVudHMsIHNvZnRlciBzaGFkb3dzLCBhbmQgbW9yZSB3aGl0ZSBzcGFjZS4KICAgLSBGb250cyBsb29rIGRhdGVkOyBjb25zaWRlciB1c2luZyBhIG1vZGVybiwgd2ViLW9wdGltaXplZCBmb250IGZvciBib2R5IHRleHQgKGUuZy4sIEludGVyLCBSb2JvdG8sIE9wZW4gU2FucykuCjMuICoqTmF2aWdhdGlvbjoqKgogICAtIFRoZSBob3Jpem9udGFsIG5hdmlnYXRpb24gYmFyIHVzZXMgYWxsLWNhcHMgdGV4dCwgd2hpY2ggY2FuIGJlIGhhcmRlciB0byByZWFkIGFuZCBsZXNzIGFjY2Vzc2libGUuIENvbnNpZGVyIHVzaW5nIHNlbnRlbmNlIGNhc2Ugb3IgdGl0bGUgY2FzZS4KICAgLSBUaGUg4oCcQ09VTlNFTElORyBSRVFVRVNU4oCdIGJ1dHRvbiBzaG91bGQgc3RhbmQgb3V0IGJ1dCBmZWVscyB2aXN1YWxseSBkaXNjb25uZWN0ZWQgZnJvbSB0aGUgbmF2aWdhdGlvbiBiYXIgZHVlIHRvIGNvbG9yIGFuZCBwbGFjZW1lbnQuIENvbnNpc3RlbnQgYnV0dG9uIHN0eWxlcyBhbmQgYmV0dGVyIHNwYWNpbmcgd291bGQgaGVscC4KNC4gKipDYXJvdXNlbC9TbGlkZXI6KioKICAgLSBUaGUgY2Fyb3VzZWwgYXJlYSBpcyBtb3N0bHkgd2hpdGUgc3BhY2UgYW5kIGxhY2tzIHZpc3VhbCBlbmdhZ2VtZW50LiBUZXh0IGlzIGhhcmQgdG8gcmVhZCBkdWUgdG8gbG93IGNvbnRyYXN0IGFnYWluc3QgdGhlIGJhY2tncm91bmQuCiAgIC0gQ2Fyb3VzZWwgYXJyb3dzIGFyZSB0b28gc21hbGwvbGlnaHQg4oCUIGluY3JlYXNlIHRoZWlyIHNpemUgYW5kIGNvbnRyYXN0LgogICAtIElmIHBvc3NpYmxlLCBhZGQgcmVsZXZhbnQgaW1hZ2VyeSBvciBiYWNrZ3JvdW5kIGdyYXBoaWNzIHRvIG1ha2UgdGhpcyBzZWN0aW9uIG1vcmUgZW5nYWdpbmcuCjUuICoqSWNvbm9ncmFwaHkgYW5kIENhcmRzOioqCiAgIC0gVGhlIGljb25zIGZvciBJbm5vdmF0ZS9MYXVuY2gvZXRjLiBhcmUgbm90IHZpc3VhbGx5IGFsaWduZWQgYW5kIHNlZW0gbWlzbWF0Y2hlZCBpbiBzdHlsZS4KICAgLSBDYXJkIFVJIGNvdWxkIGJlbmVmaXQgZnJvbSBpbmNyZWFzZWQgcGFkZGluZywgcm91bmRlZCBjb3JuZXJzLCBhbmQgbW9yZSBzZXBhcmF0aW9uIGZyb20gdGhlIGJhY2tncm91bmQgdG8gaW1wcm92ZSByZWFkYWJpbGl0eSBhbmQgdG91Y2gtZnJpZW5kbGluZXNzLgo2LiAqKlNwb25zb3IgU2VjdGlvbjoqKgogICAtIExvZ29zIGFyZSB1bmV2ZW5seSBzY2FsZWQgYW5kIG5vdCBldmVubHkgc3BhY2VkLiBVc2UgY29uc2lzdGVudCBzaXppbmcgYW5kIHBhZGRpbmcgZm9yIGEgbW9yZSBwb2xpc2hlZCBhcHBlYXJhbmNlLgogICAtIFNlY3Rpb24gdGl0bGUgKOKAnEFubnVhbCBLZXkgU3BvbnNvcnPigJ0pIGlzIG5vdCB2aXN1YWxseSBwcm9taW5lbnTigJR1c2UgbW9yZSBoaWVyYXJjaHkgdG8gbWFrZSBpdCBzdGFuZCBvdXQuCjcuICoqRm9vdGVyOioqCiAgIC0gQ29udGFjdCBpbmZvIGlzIHVzYWJsZSBidXQgY2FuIGJlIGVuaGFuY2VkIHdpdGggY2xpY2thYmxlIHBob25lIGFuZCBlbWFpbCBsaW5rcy4KICAgLSBTZWN0aW9uIGhlYWRlcnMgKOKAnENvbnRhY3QgVXMs4oCdIOKAnFJ1dGdlcnMgVW5pdmVyc2l0eeKAk0NhbWRlbuKAnSkgbmVlZCBiZXR0ZXIgdmlzdWFsIHNlcGFyYXRpb24gZnJvbSB0aGUgcmVzdCBvZiB0aGUgY29udGVudC4KICAgLSBDb25zaWRlciB1c2luZyBpY29ucyBmb3IgY29udGFjdCBtZXRob2RzLgogICAtIEZvb3RlciBiYWNrZ3JvdW5kIGNvbG9yIHdvdWxkIGhlbHAgZGlmZmVyZW50aWF0ZSBpdCBmcm9tIG1haW4gY29udGVudC4KOC4gKipBY2Nlc3NpYmlsaXR5OioqCiAgIC0gVGV4dCBjb250cmFzdCBjb3VsZCBiZSBpbXByb3ZlZCBpbiBzZXZlcmFsIGFyZWFzIChlc3BlY2lhbGx5IHRoZSBjYXJvdXNlbCkuCiAgIC0gTmF2aWdhdGlvbiBzaG91bGQgaGF2ZSB2aXNpYmxlIGZvY3VzIHN0YXRlcyBmb3Iga2V5Ym9hcmQgdXNlcnMuCiAgIC0gSW5jcmVhc2Ugb3ZlcmFsbCBmb250IHNpemUgZm9yIGJvZHkgdGV4dCBhbmQgbmF2aWdhdGlvbi4KICAgLSBFbnN1cmUgYWxsIGFjdGlvbmFibGUgZWxlbWVudHMgKGUuZy4sIGJ1dHRvbnMsIGxpbmtzKSBtZWV0IG1pbmltdW0gcmVjb21tZW5kZWQgdGFwL2NsaWNrIHNpemVzLgo5LiAqKlJlc3BvbnNpdmVuZXNzOioqCiAgIC0gV2hpbGUgc2NyZWVuc2hvdCBkb2VzbuKAmXQgc2hvdyBtb2JpbGUsIGVuc3VyZSBjYXJkcyBhbmQgbmF2YmFyIHN0YWNrIG5pY2VseSBvbiBzbWFsbGVyIHNjcmVlbnMsIGFuZCBsZWdpYmlsaXR5IGlzIHJldGFpbmVkLgoqKlN1bW1hcnk6KioKVGhlIHNpdGUgaXMgZnVuY3Rpb25hbCBidXQgdmlzdWFsbHkgb3V0ZGF0ZWQsIHdpdGggYnJhbmRpbmcgaW5jb25zaXN0ZW5jaWVzLCBsYWNrIG9mIG1vZGVybiBVSSBwb2xpc2gsIGFuZCBzZXZlcmFsIGFjY2Vzc2liaWxpdHkgc2hvcnRjb21pbmdzLiBVcGRhdGluZyB0byBjdXJyZW50IFVJIHRyZW5kcywgaW1wcm92aW5nIGNvbnNpc3RlbmN5LCBhbmQgZW5oYW5jaW5nIGFjY2Vzc2liaWxpdHkgd291bGQgcHJvdmlkZSBzaWduaWZpY2FudCBpbXByb3ZlbWVudHMu"
Set-CellB64 "O2" "KipCdXNpbmVzcyBDb250ZXh0IGFuZCBBbmFseXNpczoqKiAgClRoZSBOSiBTbWFsbCBCdXNpbmVzcyBEZXZlbG9wbWVudCBDZW50ZXIgKE5KU0JEQykgYXQgUnV0Z2Vycy1DYW1kZW4gaXMgYSBub25wcm9maXQgcmVzb3VyY2UgY2VudGVyIGRlZGljYXRlZCB0byBzdXBwb3J0aW5nIGVudHJlcHJlbmV1cnMgYW5kIHNtYWxsIGJ1c2luZXNzIG93bmVycyBhY3Jvc3MgQnVybGluZ3RvbiwgQ2FtZGVuLCBHbG91Y2VzdGVyLCBhbmQgU2FsZW0gQ291bnRpZXMgaW4gTmV3IEplcnNleS4gT3BlcmF0aW5nIHVuZGVyIFJ1dGdlcnMgVW5pdmVyc2l0eSBhbmQgcGFydGlhbGx5IGZ1bmRlZCBieSB0aGUgU0JBLCB0aGUgY2VudGVyIHByb3ZpZGVzIG5vLWNvc3QgY29uc3VsdGluZywgYnVzaW5lc3MgdHJhaW5pbmcgd29ya3Nob3BzLCBhbmQgZWR1Y2F0aW9uYWwgcmVzb3VyY2VzIGNvdmVyaW5nIGJ1c2luZXNzIHBsYW5uaW5nLCBmaW5hbmNpbmcsIG1hcmtldGluZywgY2VydGlmaWNhdGlvbnMsIGUtYnVzaW5lc3MsIGFuZCBtb3JlLiBUaGVpciBoeWJyaWQgc3VwcG9ydCBtb2RlbCBpbmNsdWRlcyBwZXJzb25hbGl6ZWQgY29uc3VsdGluZyBhbmQgZ3JvdXAgdHJhaW5pbmcsIHBvc2l0aW9uaW5nIHRoZW0gYXMgYSBrZXkgY2F0YWx5c3QgZm9yIHNtYWxsIGJ1c2luZXNzIGdyb3d0aCBhbmQgbG9jYWwgZWNvbm9taWMgZGV2ZWxvcG1lbnQuCioqUXVpY2sgV2lucyAmIFJlY29tbWVuZGVkIE9wcG9ydHVuaXRpZXM6KiogIApHaXZlbiB0aGUgY2VudGVyJ3MgZm9jdXMgb24gZWR1Y2F0aW9uIGFuZCBvdXRyZWFjaCwgaW1wbGVtZW50aW5nIGF1dG9tYXRlZCBhcHBvaW50bWVudCBzY2hlZHVsaW5nLCBDUk0gZm9yIHRyYWNraW5nIGNsaWVudCB0b3VjaHBvaW50cywgYW5kIG9ubGluZSBpbnRha2UgZm9ybXMgd291bGQgc3RyZWFtbGluZSBjbGllbnQgZW5nYWdlbWVudCBhbmQgaW1wcm92ZSBvcGVyYXRpb25hbCBlZmZpY2llbmN5LiBXZWJzaXRlIGltcHJvdmVtZW50cyBjb3VsZCBpbmNsdWRlIGEgY2xlYXJlciBjYWxsIHRvIGFjdGlvbiBmb3Igc2VydmljZSBzaWduLXVwcywgaW50ZWdyYXRpb24gb2YgbGl2ZSBjaGF0IG9yIGNoYXRib3Qgc3VwcG9ydCwgYW5kIGJldHRlciBvcmdhbml6YXRpb24gb2YgZG93bmxvYWRhYmxlIHJlc291cmNlcyBmb3IgcXVpY2sgYWNjZXNzLiBIaWdobGlnaHRpbmcgbWVhc3VyYWJsZSBjbGllbnQgc3VjY2VzcyBzdG9yaWVzIGFuZCBwcm92aWRpbmcgbW9yZSBpbnR1aXRpdmUgbmF2aWdhdGlvbiB3aWxsIGZ1cnRoZXIgZW5nYWdlIHZpc2l0b3JzIGFuZCByZWluZm9yY2UgdGhlIFNCREPigJlzIGltcGFjdC4gQWRkaXRpb25hbGx5LCByZXZpZXdpbmcgYmFjay1vZmZpY2Ugd29ya2Zsb3dzIGZvciBhdXRvbWF0aW9uIChlLmcuLCB0cmFpbmluZyByZWdpc3RyYXRpb24sIGV2ZW50IHJlbWluZGVycywgYW5kIGZlZWRiYWNrIGNvbGxlY3Rpb24pIGNvdWxkIHNhdmUgc3RhZmYgdGltZSBhbmQgZGVsaXZlciBhIG1vcmUgc2VhbWxlc3MgZXhwZXJpZW5jZSB0byBzbWFsbCBidXNpbmVzcyBjbGllbnRzLg=="
Set-CellB64 "Q2" "U3ViamVjdDogQm9vc3RpbmcgSW1wYWN0ICYgRWZmaWNpZW5jeSBmb3IgTkpTQkRDIGF0IFJ1dGdlcnMtQ2FtZGVuCkhpIFtOYW1lXSwKSeKAmW0gaW1wcmVzc2VkIGJ5IGhvdyBOSlNCREMgQCBSdXRnZXJzLUNhbWRlbiBlbXBvd2VycyBTb3V0aCBKZXJzZXkgZW50cmVwcmVuZXVycyB3aXRoIGV4cGVydCBzdXBwb3J0IGFuZCByZXNvdXJjZXMuIFdlIGhlbHAgbm9ucHJvZml0IGNlbnRlcnMgbGlrZSB5b3VycyBzdHJlYW1saW5lIGNsaWVudCBlbmdhZ2VtZW50IGJ5IGF1dG9tYXRpbmcgYXBwb2ludG1lbnQgc2NoZWR1bGluZyBhbmQgaW50YWtlIGZvcm1zLCByZWRlc2lnbmluZyB3ZWJzaXRlcyBmb3IgZWFzaWVyIHJlc291cmNlIGFjY2VzcywgYW5kIGF1dG9tYXRpbmcgZXZlbnQgZm9sbG93LXVwc+KAlGZyZWVpbmcgdXAgc3RhZmYgZm9yIGhpZ2hlci1pbXBhY3Qgd29yay4KT3BlbiB0byBhIHF1aWNrIGNhbGwgdG8gZXhwbG9yZSBob3cgc21hbGwgYXV0b21hdGlvbiB3aW5zIGNhbiBtYWtlIGEgYmlnIGRpZmZlcmVuY2UgZm9yIHlvdXIgdGVhbSBhbmQgdGhlIGJ1c2luZXNzIG93bmVycyB5b3Ugc2VydmU/"
$ws.Range("E2").Value = 5
$ws.Range("F2").Value = 7
$ws.Range("M2").Value = 2

# ----- Row 3 -----
Set-CellB64 "A3" "Q2hJSlZRMXcyTXI5d29rUnJ0Mm5ZcVhHNG8w"
Set-CellB64 "B3" "TkpTQkRDIGF0IFdpbGxpYW0gUGF0ZXJzb24gVW5pdmVyc2l0eQ=="
Set-CellB64 "C3" "Y29uc3VsdGFudCwgcG9pbnRfb2ZfaW50ZXJlc3QsIGVzdGFibGlzaG1lbnQ="
Set-CellB64 "D3" "KDk3MykgMzIxLTEzNzg="
Set-CellB64 "H3" "aHR0cHM6Ly9zYmRjd3B1LmNvbS8="
Set-CellB64 "I3" "aHR0cHM6Ly9tYXBzLmdvb2dsZS5jb20vP2NpZD0xMDIyMzk1MjUxNzcxMTUxOTE1MCZnX21wPUNpZG5iMjluYkdVdWJXRndjeTV3YkdGalpYTXVkakV1VUd4aFkyVnpMbE5sWVhKamFGUmxlSFFRQWhnRUlBQQ=="
Set-CellB64 "K3" "SSBjYW5ub3QgZXhwcmVzcyBlbm91Z2ggZ3JhdGl0dWRlIGZvciBEYXluYSBDYXB1dG8gYW5kIGhlciBvZmZpY2UgZm9yIHRoZWlyIGluY3JlZGlibGUgc3VwcG9ydCBpbiBoZWxwaW5nIG1lIGNvbXBsZXRlIG15IGdyYW50IGFwcGxpY2F0aW9uIGZvciBNZXJjZWQgTWFub3IuIEFzIGEgc21hbGwsIHdvbWVuLW93bmVkIGJ1c2luZXNzLCBuYXZpZ2F0aW5nIHRoZSBncmFudCBwcm9jZXNzIGNhbiBiZSBvdmVyd2hlbG1pbmcsIGJ1dCBEYXluYeKAmXMgcGF0aWVuY2UsIGV4cGVydGlzZSwgYW5kIHVud2F2ZXJpbmcgc3VwcG9ydCBtYWRlIGFsbCB0aGUgZGlmZmVyZW5jZS4gU2hlIGd1aWRlZCBtZSBldmVyeSBzdGVwIG9mIHRoZSB3YXksIGVuc3VyaW5nIEkgaGFkIHRoZSByZXNvdXJjZXMgYW5kIGtub3dsZWRnZSBuZWVkZWQgdG8gc3VjY2VlZC4gVGhpcyBraW5kIG9mIHN1cHBvcnQgaXMgdHJlbWVuZG91cyBmb3Igc21hbGwgYnVzaW5lc3NlcyBsaWtlIG1pbmUsIHByb3ZpZGluZyBvcHBvcnR1bml0aWVzIHRoYXQgY2FuIHRydWx5IGJlIHRyYW5zZm9ybWF0aXZlLiBUaGFuayB5b3UsIERheW5hLCBmb3IgeW91ciBkZWRpY2F0aW9uIGFuZCBmb3IgY2hhbXBpb25pbmcgc21hbGwgYnVzaW5lc3Mgb3duZXJz4oCUeW91IGFyZSBtYWtpbmcgYW4gaW52YWx1YWJsZSBpbXBhY3QhIHwgSSBhbSBzbyBob25vcmVkIHRvIGhhdmUgYmVlbiBhIHBhcnQgb2YgdGhpcyBhbWF6aW5nIG9yZ2FuaXphdGlvbi4gUmVpbmEgd2hvIGlzIG9uZSBvZiB0aGUgbWVudG9ycyBoZWxwZWQgbWUgYnkgZ2V0dGluZyBteSBjZXJ0aWZpY2F0aW9ucyBhbmQgbXkgR29vZ2xlIGJ1c2luZXNzIHVwIGFuZCBydW5uaW5nLiBJIHdhcyB3aXRoIHRoZSBvcmdhbml6YXRpb24gZm9yIGxlc3MgdGhhbiBhIHllYXIgYW5kIHdvbiBhIFNCREMgQXdhcmQsIEZlZEV4IGdyYW50IGFuZCBnb3QgYWNjZXB0ZWQgaW50byBUYXJnZXQuIEZvciBhbnlvbmUgd2hvIGlzIGxvb2tpbmcgZm9yIGNvbW11bml0eSBhbmQgc3VwcG9ydCwgc2lnbiB1cCB0b2RheSEgVGhleSB3ZXJlIGFsbCBncmVhdCEgVGhhbmsgeW91IG9uY2UgYWdhaW4gdG8gZXZlcnlvbmUhIHwgRXh0cmVtZWx5IGhlbHBmdWwgZm9yIG5ldyBzdGFydCB1cCBidXNpbmVzc2VzLiAgVmVyeSBpbmZvcm1hdGl2ZS4="
Set-CellB64 "L3" "c2JkY0B3cHVuai5lZHU="
Set-CellB64 "N3" "KipVSSBJbXByb3ZlbWVudHM6KioKMS4gKipBYm92ZS10aGUtRm9sZCBDb250ZW50Kio6CiAgIC0gVGhlIHRvcCBwb3J0aW9uIG9mIHRoZSBzaXRlIGlzIG1vc3RseSBlbXB0eSBzcGFjZSwgd2hpY2ggY3JlYXRlcyBhIHdlYWsgZmlyc3QgaW1wcmVzc2lvbiBhbmQgZm9yY2VzIHVzZXJzIHRvIHNjcm9sbC4gQ29uc2lkZXIgZmlsbGluZyB0aGUgaGVybyBzZWN0aW9uIHdpdGggYSBjbGVhciBjYWxsLXRvLWFjdGlvbiwgYW4gZW5nYWdpbmcgaGVhZGluZywgb3IgYSBoZXJvIGltYWdlL2dyYXBoaWMgdGhhdCBjb21tdW5pY2F0ZXMgeW91ciBicmFuZCBtZXNzYWdlLgoyLiAqKkJyYW5kaW5nICYgQ29oZXNpdmVuZXNzKio6CiAgIC0gV2hpbGUgdGhlIGNvbG9yIHBhbGV0dGUgaXMgY29uc2lzdGVudCwgdGhlIHBpbmsgb3ZlcmxheSBvbiBpbWFnZXMgZG9lcyBub3QgYWxpZ24gd2VsbCB3aXRoIHRoZSByZXN0IG9mIHRoZSBicmFuZOKAmXMgcHJvZmVzc2lvbmFsIGFuZCBjbGVhbiBuYXZ5L3doaXRlIGxvb2suIFRyeSBzdWJ0bGVyIGltYWdlIG92ZXJsYXlzIG9yIHVzZSByZWFsLCB2aWJyYW50IHBob3Rvcy4KMy4gKipOYXZpZ2F0aW9uL01lbnUgQmFyKio6CiAgIC0gVGhlIG5hdmlnYXRpb24gYmFyIGZvbnQgYXBwZWFycyBzbWFsbCBhbmQgc3BhY2VkIGZhciBhcGFydCwgd2hpY2ggbWF5IG1ha2UgcmVhZGluZyBhbmQgdGFyZ2V0aW5nIGxpbmtzIGRpZmZpY3VsdCwgZXNwZWNpYWxseSBvbiBsYXJnZXIgc2NyZWVucy4gSW5jcmVhc2UgZm9udCBzaXplIGFuZCBwYWRkaW5nIGZvciBpbXByb3ZlZCB1c2FiaWxpdHkuCiAgIC0gVGhlIHNvY2lhbCBtZWRpYSBpY29ucyBhcmUgbWluaW1hbGx5IHZpc2libGU7IGluY3JlYXNpbmcgdGhlaXIgc2l6ZSBvciBjaGFuZ2luZyBjb2xvciBvbiBob3ZlciB3b3VsZCBpbXByb3ZlIHZpc2liaWxpdHkuCjQuICoqQ29udGVudCBEZW5zaXR5IGFuZCBCYWxhbmNlKio6CiAgIC0gVGhlIGxvbmcgc2Nyb2xsIGJldHdlZW4gdGhlIGhlYWRlciBhbmQgYWN0dWFsIGNvbnRlbnQgc2VjdGlvbnMgKGUuZy4sICJSZWNlbnQgRXZlbnRzIikgZmVlbHMgdW5iYWxhbmNlZC4gVGlnaHRlbiB1cCB2ZXJ0aWNhbCBzcGFjaW5nIHRvIG1ha2UgYmV0dGVyIHVzZSBvZiBzY3JlZW4gcmVhbCBlc3RhdGUuCiAgIC0gVGhlIOKAnFJlY2VudCBFdmVudHPigJ0gbW9kdWxlIGFwcGVhcnMgd2l0aCBhIGJsYW5rIHNwYWNlIGJlbG93IGl0LCB3aGljaCBpcyBjb25mdXNpbmcuIEVpdGhlciBwb3B1bGF0ZSB0aGlzIGFyZWEgb3IgcmVtb3ZlIHRoZSBoZWFkaW5nIHVudGlsIGV2ZW50cyBleGlzdC4KNS4gKipBY2Nlc3NpYmlsaXR5Kio6CiAgIC0gV2hpbGUgdGhlcmXigJlzIGFuIGFjY2Vzc2liaWxpdHkgcGx1Z2luL3dpZGdldCB2aXNpYmxlLCBzb21lIGNvbG9yIGNvbnRyYXN0IChsaWtlIHdoaXRlIHRleHQgb24gcmVkIG9yIHllbGxvdyBzdGFycyBvbiBuYXZ5KSBjb3VsZCBiZSBpbXByb3ZlZOKAlGVzcGVjaWFsbHkgZm9yIHZpc3VhbGx5IGltcGFpcmVkIHVzZXJzLgogICAtIENvbnNpZGVyIGVuc3VyaW5nIGJ1dHRvbnMgYW5kIGxpbmtzIGFyZSBrZXlib2FyZC90YWIgYWNjZXNzaWJsZSB3aXRoIHZpc2libGUgZm9jdXMgaW5kaWNhdG9ycy4KNi4gKipNb2Rlcm4gVUkgVHJlbmRzKio6CiAgIC0gVGhlIGxhcmdlIGVtcHR5IGFyZWFzIGFuZCBibG9ja3Mgb2YgY29sb3IgZmVlbCBvdXRkYXRlZDsgbW9kZXJuIFVJcyB1c2UgbW9yZSB3aGl0ZSBzcGFjZSBzdHJhdGVnaWNhbGx5LCBsYXllcmluZywgc2xpZ2h0IHNoYWRvd3MsIGFuZCBncmFkaWVudHMgZm9yIGRlcHRoLgogICAtIENhcmQgZWxlbWVudHMgYXQgdGhlIGJvdHRvbSBjb3VsZCB1c2Ugc29mdCBkcm9wIHNoYWRvd3Mgb3Igc2xpZ2h0bHkgcm91bmRlZCBjb3JuZXJzIGZvciBhIG1vcmUgY29udGVtcG9yYXJ5IGxvb2suCiAgIC0gQ29uc2lkZXIgdXNpbmcgbW9yZSB2aXN1YWxseSBlbmdhZ2luZyBjb21wb25lbnRzIGxpa2UgaWNvbnMsIGlsbHVzdHJhdGlvbnMsIG9yIGFuaW1hdGlvbiBmb3Iga2V5IGludGVyYWN0aW9ucy4KNy4gKipGb3JtIERlc2lnbioqOgogICAtIE5ld3NsZXR0ZXIgc2lnbi11cCBmaWVsZHMgKGZpcnN0L2xhc3QgbmFtZSwgZW1haWwpIGxhY2sgY2xlYXIgbGFiZWxpbmcgYW5kIGlubGluZSB2YWxpZGF0aW9uLiBUaGV5IGNvdWxkIGJlIGFycmFuZ2VkIGhvcml6b250YWxseSBvbiB3aWRlciBzY3JlZW5zIGZvciBhIG1vcmUgY29tcGFjdCBsb29rLgogICAtIFRoZSDigJxTdWJtaXTigJ0gYnV0dG9uIGRlc2lnbiBpcyB2ZXJ5IGJhc2ljIGFuZCBlYXNpbHkgbWlzc2Vk4oCUdXNlIGEgbW9yZSBwcm9taW5lbnQgYnV0dG9uIHdpdGggaGlnaGVyIGNvbnRyYXN0IG9yIGluY2x1ZGUgYW4gaWNvbi4KOC4gKipJbWFnZSAmIENvbnRlbnQgT3ZlcmxheXMqKjoKICAgLSBUaGUgY29sb3JlZCBvdmVybGF5cyBvbiBpbWFnZXMgKGRlZXAgcmVkIGFuZCBibHVlKSBjYW4gb2JzY3VyZSBwaG90byBkZXRhaWxzLCBtYWtpbmcgdGhlIHZpc3VhbHMgbGVzcyBlbmdhZ2luZy4gVXNlIG1vcmUgc3VidGxlIG92ZXJsYXlzIGFuZC9vciBpbmNyZWFzZSBwaG90byBicmlnaHRuZXNzLgogICAtIFNwb25zb3Igc2VjdGlvbiBsb29rcyBjb2hlc2l2ZSBidXQgY291bGQgdXNlIG1vcmUgYnJlYXRoaW5nIHJvb20gYmV0d2VlbiBsb2dvcyB0byBhdm9pZCBhIGNyYW1wZWQgZmVlbC4KOS4gKipNb2JpbGUgUmVzcG9uc2l2ZW5lc3MqKjoKICAgLSBXaGlsZSByZXNwb25zaXZlbmVzcyBjYW4ndCBiZSBjb25maXJtZWQgaW4gYSBzY3JlZW5zaG90LCBsYXJnZSBmaXhlZC13aWR0aCBzZWN0aW9ucyBhbmQgbG90cyBvZiBob3Jpem9udGFsIGNhcmRzIGNhbiBjcmVhdGUgcHJvYmxlbXMgb24gbW9iaWxlLiBNYWtlIHN1cmUgYWxsIGNhcmRzIHN0YWNrIHZlcnRpY2FsbHkgYW5kIHBhZGRpbmcvbWFyZ2lucyBhcmUgbW9iaWxlLWZyaWVuZGx5LgoxMC4gKipIaWVyYXJjaHkgJiBUeXBvZ3JhcGh5Kio6CiAgICAtIFVzZSBtb3JlIHByb25vdW5jZWQgaGVhZGluZ3MgYW5kIHN1YmhlYWRpbmdzIHRvIGJyZWFrIHVwIHNlY3Rpb25zIGFuZCBndWlkZSB0aGUgdXNlcidzIGV5ZS4KICAgIC0gRW5zdXJlIGNvbnNpc3RlbnQgaGVhZGluZyBzaXplcyBhbmQgc3BhY2luZyBmb3IgYSBwb2xpc2hlZCBhcHBlYXJhbmNlLgotLS0KKipTdW1tYXJ5OioqICAKVGhlIHNpdGUgY29tbXVuaWNhdGVzIHRydXN0IGJ1dCBmZWVscyBzbGlnaHRseSBvdXRkYXRlZCBhbmQgbGFja3MgdmlzdWFsIGVuZ2FnZW1lbnQgYWJvdmUgdGhlIGZvbGQuIFRpZ2h0ZW4gc3BhY2luZywgbW9kZXJuaXplIHZpc3VhbHMgKHNoYWRvd3MsIGdyYWRpZW50cywgY2FyZHMpLCBpbXByb3ZlIGFib3ZlLXRoZS1mb2xkIGNvbnRlbnQsIHJldmlzaXQgY29sb3IvYWNjZXNzaWJpbGl0eSBjb250cmFzdCwgYW5kIHVzZSByZWFsIGltYWdlcnkuIE1ha2UgbmF2aWdhdGlvbiwgZm9ybXMsIGFuZCBDVEFzIG1vcmUgdXNhYmxlIGFuZCB2aXNpYmxlLiBBZGp1c3Qgb3ZlcmxheXMgZm9yIGNsZWFuZXIsIG1vcmUgaW52aXRpbmcgYnJhbmRpbmcuIAoqKihOb3RlOiBBbGwgc3VnZ2VzdGlvbnMgYXJlIFVJLWZvY3VzZWQgYW5kIG5vdCBhYm91dCBjb250ZW50L2Z1bmN0aW9uYWxpdHkuKSoq"
Set-CellB64 "O3" "IyMjIDEuIEJ1c2luZXNzIElkZW50aXR5LCBJbmR1c3RyeSwgYW5kIFBvc2l0aW9uaW5nCioqQnVzaW5lc3M6KiogIApUaGUgTmV3IEplcnNleSBTbWFsbCBCdXNpbmVzcyBEZXZlbG9wbWVudCBDZW50ZXIgKE5KU0JEQykgYXQgV2lsbGlhbSBQYXRlcnNvbiBVbml2ZXJzaXR5CioqSW5kdXN0cnk6KiogIApCdXNpbmVzcyBhbmQgZWNvbm9taWMgZGV2ZWxvcG1lbnQsIGZvY3VzaW5nIG9uIHNtYWxsIGJ1c2luZXNzIHN1cHBvcnQKKipQb3NpdGlvbmluZzoqKiAgCkEgcHVibGljbHkgYW5kIHByaXZhdGVseSBmdW5kZWQgcmVzb3VyY2UgaHViIG9mZmVyaW5nIG5vLWNvc3QgY291bnNlbGluZywgZWR1Y2F0aW9uLCB0cmFpbmluZywgYW5kIHByb2Zlc3Npb25hbCBhZHZpY2UgdG8gc21hbGwgYnVzaW5lc3NlcyBhbmQgZW50cmVwcmVuZXVycyBpbiBQYXNzYWljIGFuZCBVbmlvbiBjb3VudGllcywgTmV3IEplcnNleS4gVGhlIGNlbnRlciBvcGVyYXRlcyB3aXRoaW4gV2lsbGlhbSBQYXRlcnNvbiBVbml2ZXJzaXR54oCZcyBDb3RzYWtvcyBDb2xsZWdlIG9mIEJ1c2luZXNzIGFuZCBwYXJ0bmVycyB3aXRoIHRoZSBTQkEgYW5kIG90aGVyIG9yZ2FuaXphdGlvbnMgdG8gaGVscCBsb2NhbCBidXNpbmVzc2VzIHN0YXJ0LCBncm93LCBhbmQgdGhyaXZlLgotLS0KIyMjIDIuIFBlcnNvbmFsaXplZCBCcmllZiB3aXRoIE9wcG9ydHVuaXRpZXMKVGhlIE5KU0JEQyBhdCBXaWxsaWFtIFBhdGVyc29uIFVuaXZlcnNpdHkgcGxheXMgYSBjcml0aWNhbCByb2xlIGluIGVtcG93ZXJpbmcgc21hbGwgYnVzaW5lc3Mgb3duZXJzIGFuZCBlbnRyZXByZW5ldXJzIGFjcm9zcyBub3J0aGVybiBOZXcgSmVyc2V5IHdpdGggbm8tY29zdCBleHBlcnQgY291bnNlbGluZywgYnVzaW5lc3MgcGxhbiBkZXZlbG9wbWVudCwgYWNjZXNzIHRvIGNhcGl0YWwsIG5ldHdvcmtpbmcsIGFuZCB0cmFpbmluZyByZXNvdXJjZXMuIEdpdmVuIHRoZSBjZW50ZXLigJlzIHdpZGUtcmFuZ2luZyBzZXJ2aWNlIG9mZmVyaW5ncyBhbmQgaGlnaCBjb21tdW5pdHkgZW5nYWdlbWVudCwgdGhlcmUgYXJlIHF1aWNrIHdpbnMgYXZhaWxhYmxlIHRocm91Z2ggZW5oYW5jZWQgYXV0b21hdGlvbuKAlHN1Y2ggYXMgc3RyZWFtbGluaW5nIGFwcG9pbnRtZW50IHNjaGVkdWxpbmcsIGludGFrZSwgYW5kIGZvbGxvdy11cCBwcm9jZXNzZXMgKGUuZy4sIG9ubGluZSBib29raW5nLCBhdXRvbWF0ZWQgZW1haWwgcmVtaW5kZXJzLCBhbmQgZGlnaXRhbCByZXNvdXJjZSBkZWxpdmVyeSkuIE1vZGVybml6aW5nIHRoZSB3ZWJzaXRlIGZvciBhIG1vcmUgc2VhbWxlc3MgbmF2aWdhdGlvbiBhbmQgYWNjZXNzaWJpbGl0eSBleHBlcmllbmNlLCBhbmQgaW50ZWdyYXRpbmcgc2VsZi1zZXJ2aWNlIHRvb2xzIChzdWNoIGFzIHNtYXJ0IEZBUXMgb3IgcmVzb3VyY2UgbGlicmFyaWVzKSwgY291bGQgYWxzbyBpbXByb3ZlIGNsaWVudCBlbmdhZ2VtZW50IGFuZCByZWR1Y2UgbWFudWFsIHdvcmtsb2FkLiBBZGRpdGlvbmFsbHksIG1hcHBpbmcgd29ya2Zsb3dzIGZvciBldmVudCBwcm9tb3Rpb24sIHNpZ24tdXBzLCBhbmQgZmVlZGJhY2sgY29sbGVjdGlvbiBjb3VsZCBmdXJ0aGVyIG9wdGltaXplIGludGVybmFsIHByb2Nlc3Nlc+KAlGFtcGxpZnlpbmcgdGhlIGNlbnRlcuKAmXMgYWxyZWFkeSBzdHJvbmcgaW1wYWN0IGFuZCBmcmVlaW5nIHVwIHRpbWUgZm9yIGRlZXBlciBvbmUtb24tb25lIGFkdmlzaW5nIGFuZCBwYXJ0bmVyc2hpcCBkZXZlbG9wbWVudC4="
Set-CellB64 "P3" "KipCdXNpbmVzcyBQYWluIFBvaW50IFJlcG9ydCoqICAKKkNsaWVudDogVGhlIE5ldyBKZXJzZXkgU21hbGwgQnVzaW5lc3MgRGV2ZWxvcG1lbnQgQ2VudGVyIChOSlNCREMpIGF0IFdpbGxpYW0gUGF0ZXJzb24gVW5pdmVyc2l0eSoKLS0tCioqUGFpbiBQb2ludCAxOiBPdXRkYXRlZCBXZWJzaXRlIERlc2lnbiAmIFN1Ym9wdGltYWwgVXNlciBFeHBlcmllbmNlKiogIAogIC0gKipJbXBhY3Q6KiogVGhlIHdlYnNpdGXigJlzIGN1cnJlbnQgZGVzaWdu4oCUd2l0aCBlbXB0eSBhYm92ZS10aGUtZm9sZCBzcGFjZSwgb3V0ZGF0ZWQgdmlzdWFsIGVsZW1lbnRzLCBhbmQgd2VhayBDVEFz4oCUbWF5IGNyZWF0ZSBhIHBvb3IgZmlyc3QgaW1wcmVzc2lvbiBmb3IgbmV3IHZpc2l0b3JzLiBUaGlzIGNhbiBsZWFkIHRvIGxvc3QgZW5nYWdlbWVudCwgbG93ZXIgdHJ1c3QsIGFuZCByZWR1Y2VkIGlucXVpcmllcyBvciBldmVudCBzaWduLXVwcywgZGFtYWdpbmcgY3JlZGliaWxpdHkgYXMgYSBsZWFkaW5nIGJ1c2luZXNzIHJlc291cmNlLgogIC0gKipIb3cgV2UgQ2FuIEhlbHA6KiogKldlYnNpdGUgUmVkZXNpZ24gU2VydmljZXMqIOKAkyBXZSB3aWxsIG1vZGVybml6ZSB5b3VyIHdlYnNpdGUgZm9yIGludHVpdGl2ZSBuYXZpZ2F0aW9uLCBlbmdhZ2luZyB2aXN1YWxzLCBpbXByb3ZlZCBDVEFzLCBhbmQgY29uc2lzdGVudCBicmFuZGluZy4gRW5oYW5jZWQgYWNjZXNzaWJpbGl0eSBhbmQgcmVzcG9uc2l2ZSBkZXNpZ24gd2lsbCBlbnN1cmUgYWxsIHVzZXJzLCBpbmNsdWRpbmcgdGhvc2Ugd2l0aCBkaXNhYmlsaXRpZXMgb3Igb24gbW9iaWxlIGRldmljZXMsIGNhbiBlYXNpbHkgYWNjZXNzIHlvdXIgcmVzb3VyY2VzIGFuZCBjb25uZWN0IHdpdGggeW91ciB0ZWFtLgotLS0KKipQYWluIFBvaW50IDI6IE1hbnVhbCBhbmQgSW5lZmZpY2llbnQgT3BlcmF0aW9uYWwgUHJvY2Vzc2VzKiogIAogIC0gKipJbXBhY3Q6KiogSGlnaCBjbGllbnQgdm9sdW1lIGFuZCBkaXZlcnNlIHNlcnZpY2Ugb2ZmZXJpbmdzIGxpa2VseSByZWx5IG9uIG1hbnVhbCBzdGVwcyBmb3Igc2NoZWR1bGluZywgaW50YWtlLCBhcHBvaW50bWVudCByZW1pbmRlcnMsIGV2ZW50IHNpZ24tdXBzLCBhbmQgZm9sbG93LXVwLiBUaGlzIHdhc3RlcyBzdGFmZiB0aW1lLCBkZWxheXMgcmVzcG9uc2VzLCBpbmNyZWFzZXMgdGhlIHJpc2sgb2YgZHJvcHBlZCBsZWFkcywgYW5kIGxpbWl0cyB0aGUgY2VudGVy4oCZcyBhYmlsaXR5IHRvIHNjYWxlIHN1cHBvcnQgb3IgaW5jcmVhc2UgY291bnNlbGluZyBjYXBhY2l0eS4KICAtICoqSG93IFdlIENhbiBIZWxwOioqICpQcm9jZXNzIEF1dG9tYXRpb24gKE9wZXJhdGlvbnMgJiBTdXBwb3J0KSog4oCTIFdlIHdpbGwgaW1wbGVtZW50IG9ubGluZSBhcHBvaW50bWVudCBib29raW5nLCBhdXRvbWF0ZWQgaW50YWtlIHF1ZXN0aW9ubmFpcmVzLCB0cmlnZ2VyZWQgZW1haWwgcmVtaW5kZXJzLCBzZWxmLXNlcnZpY2UgcmVzb3VyY2UgZGVsaXZlcnksIGFuZCBhdXRvbWF0ZWQgZmVlZGJhY2sgY29sbGVjdGlvbi4gVGhpcyB3aWxsIHJlZHVjZSBhZG1pbmlzdHJhdGl2ZSB3b3JrbG9hZCBhbmQgZW5hYmxlIHlvdXIgY29uc3VsdGFudHMgdG8gaW52ZXN0IG1vcmUgdGltZSBpbiBoaWdoLXZhbHVlIGFkdmlzaW5nLgotLS0KKipQYWluIFBvaW50IDM6IExpbWl0ZWQgU2VsZi1TZXJ2aWNlIERpZ2l0YWwgVG9vbHMgJiBLbm93bGVkZ2UgQWNjZXNzKiogIAogIC0gKipJbXBhY3Q6KiogVGhlIGN1cnJlbnQgd2Vic2l0ZSBvZmZlcnMgbGltaXRlZCBpbnRlcmFjdGl2ZSBmZWF0dXJlcyAoZS5nLiwgbm8gc21hcnQgRkFRLCByZXNvdXJjZSBsaWJyYXJpZXMsIG9yIGxpdmUgY2hhdCkuIFRoaXMgbWVhbnMgcHJvc3BlY3RpdmUgY2xpZW50cyB3aXRoIGNvbW1vbiBxdWVzdGlvbnMgbXVzdCB3YWl0IGZvciBhIHN0YWZmIHJlc3BvbnNlLCBjYXVzaW5nIGRlbGF5cywgaW5jcmVhc2luZyBzdGFmZiB3b3JrbG9hZCwgYW5kIHBvdGVudGlhbGx5IGRldGVycmluZyBuZXcgc2lnbi11cHMgb3IgZm9sbG93LXVwcy4KICAtICoqSG93IFdlIENhbiBIZWxwOioqICpEaWdpdGFsIFNlbGYtU2VydmljZSBJbnRlZ3JhdGlvbiAmIENvbnN1bHRpbmcqIOKAkyBXZeKAmWxsIGVxdWlwIHlvdXIgc2l0ZSB3aXRoIHNtYXJ0IGtub3dsZWRnZSBiYXNlcywgZHluYW1pYyBGQVFzLCBhbmQgYW4gYWNjZXNzaWJsZSByZXNvdXJjZSBwb3J0YWwsIHBsdXMgY2hhdGJvdHMgZm9yIGNvbW1vbiBpbnF1aXJpZXMuIFRoaXMgZW1wb3dlcnMgY2xpZW50cyB0byBmaW5kIGluZm9ybWF0aW9uIGluc3RhbnRseSBhbmQgYWxsb3dzIHN0YWZmIHRvIGZvY3VzIG9uIGNvbXBsZXggb3IgaGlnaC10b3VjaCBzdXBwb3J0LgotLS0KKipQYWluIFBvaW50IDQ6IExhY2sgb2YgV29ya2Zsb3cgTWFwcGluZyBmb3IgRXZlbnQgTWFuYWdlbWVudCAmIEZlZWRiYWNrKiogIAogIC0gKipJbXBhY3Q6KiogRXZlbnQgcHJvbW90aW9uLCByZWdpc3RyYXRpb24sIGFuZCBmZWVkYmFjayBjb2xsZWN0aW9uIG1heSBiZSBoYW5kbGVkIG1hbnVhbGx5IG9yIHdpdGggYWQtaG9jIHRvb2xzLCBsZWFkaW5nIHRvIGluY29uc2lzdGVudCBvdXRyZWFjaCwgbG9zdCByZWdpc3RyYXRpb25zLCBwb3RlbnRpYWxseSBsb3dlciBhdHRlbmRhbmNlLCBhbmQgbWlzc2VkIGluc2lnaHRzIG9uIGV2ZW50IGVmZmVjdGl2ZW5lc3MuCiAgLSAqKkhvdyBXZSBDYW4gSGVscDoqKiAqT3JnYW5pemF0aW9uYWwvUHJvY2VzcyBDb25zdWx0aW5nICYgRXZlbnQgQXV0b21hdGlvbiog4oCTIFdlIGNhbiBtYXAsIGRvY3VtZW50LCBhbmQgc3RyZWFtbGluZSB5b3VyIGV2ZW50IHdvcmtmbG93cywgaW50ZWdyYXRpbmcgYXV0b21hdGVkIHByb21vdGlvbiwgc2lnbi11cCBmb3JtcywgYW5kIHBvc3QtZXZlbnQgZmVlZGJhY2sgc3VydmV5cy4gVGhpcyBib29zdHMgZXZlbnQgc3VjY2VzcywgcmVkdWNlcyBlZmZvcnQsIGFuZCBkZWxpdmVycyBhY3Rpb25hYmxlIGFuYWx5dGljcyBmb3Igb25nb2luZyBpbXByb3ZlbWVudC4KLS0tCioqUGFpbiBQb2ludCA1OiBQb3RlbnRpYWwgR2FwcyBpbiBGaW5hbmNpYWwgUHJvY2VzcyBFZmZpY2llbmN5KiogIAogIC0gKipJbXBhY3Q6KiogV2hpbGUgbm90IHZpc2libGUgaW4gdGhlIFVJIG9yIHJldmlld3MsIHRoZSBjZW50ZXIgc3VwcG9ydHMgc21hbGwgYnVzaW5lc3MgZmluYW5jZXMgYW5kIG1heSBpdHNlbGYgYmVuZWZpdCBmcm9tIHN0cmVhbWxpbmVkLCB0cmFuc3BhcmVudCBib29ra2VlcGluZyBhbmQgcmVwb3J0aW5n4oCUZXNwZWNpYWxseSBhcyBwYXJ0bmVyc2hpcHMgYW5kIGZ1bmRpbmcgc3RyZWFtcyBncm93LiBJbmVmZmljaWVudCBwcm9jZXNzZXMgaGVyZSBjb3VsZCBsZWFkIHRvIG1pc3NlZCBncmFudC9wYXJ0bmVyc2hpcCBvcHBvcnR1bml0aWVzIG9yIHJlcG9ydGluZyBib3R0bGVuZWNrcywgcmVkdWNpbmcgb3JnYW5pemF0aW9uYWwgYWdpbGl0eS4KICAtICoqSG93IFdlIENhbiBIZWxwOioqICpCb29ra2VlcGluZyAmIEFjY291bnRpbmcgU2VydmljZXMqIOKAkyBXZeKAmWxsIGFzc2VzcyBhbmQgb3B0aW1pemUgeW91ciBpbnRlcm5hbCBmaW5hbmNpYWwgdHJhY2tpbmcgYW5kIHJlcG9ydGluZywgdXNpbmcgY2xvdWQtYmFzZWQgdG9vbHMgZm9yIHJlYWwtdGltZSBpbnNpZ2h0cywgY2xlYW5lciBhdWRpdHMsIGFuZCBlYXNpZXIgZ3JhbnQgY29tcGxpYW5jZS4gVGhpcyBlbmhhbmNlcyBvcmdhbml6YXRpb25hbCBwcm9mZXNzaW9uYWxpc20gYW5kIHBvc2l0aW9ucyBOSlNCREMgYXMgYSBtb2RlbCBmb3IgdGhlIGNsaWVudHMgaXQgc3VwcG9ydHMuCi0tLQoqKlN1bW1hcnkqKiAgCllvdXIgdGVhbeKAmXMgcmVwdXRhdGlvbiBmb3Igc3VwcG9ydGl2ZSwgdHJhbnNmb3JtYXRpdmUgZ3VpZGFuY2UgaXMgY2xlYXIgZnJvbSBzdGVsbGFyIHJldmlld3MuIE1vZGVybml6aW5nIHlvdXIgZGlnaXRhbCBwcmVzZW5jZSBhbmQgYXV0b21hdGluZyBrZXkgdG91Y2hwb2ludHMgd2lsbCBlbnN1cmUgTkpTQkRD4oCZcyBvcGVyYXRpb25hbCBleGNlbGxlbmNlIG1hdGNoZXMgdGhlIHBhc3Npb24gb2YgaXRzIHN0YWZm4oCUZGVsaXZlcmluZyBldmVuIGdyZWF0ZXIgaW1wYWN0IHRvIE5ldyBKZXJzZXnigJlzIHNtYWxsIGJ1c2luZXNzIGNvbW11bml0eS4gV2UgbG9vayBmb3J3YXJkIHRvIHN1cHBvcnRpbmcgeW91ciBuZXh0IHBoYXNlIG9mIGdyb3d0aC4="
Set-CellB64 "Q3" "U3ViamVjdDogRWxldmF0ZSBOSlNCREPigJlzIEltcGFjdCB3aXRoIFNlYW1sZXNzIERpZ2l0YWwgU29sdXRpb25zCkhpIFtOYW1lXSwKSeKAmW0gaW1wcmVzc2VkIGJ5IE5KU0JEQyBhdCBXaWxsaWFtIFBhdGVyc29u4oCZcyBkZWRpY2F0aW9uIHRvIGVtcG93ZXJpbmcgc21hbGwgYnVzaW5lc3NlcyBpbiBQYXNzYWljIGFuZCBVbmlvbiBjb3VudGllcy4gV2UgY2FuIGhlbHAgYW1wbGlmeSB5b3VyIGltcGFjdCBieSBtb2Rlcm5pemluZyB5b3VyIHdlYnNpdGUgZm9yIGJldHRlciBhY2Nlc3NpYmlsaXR5IGFuZCBlbmdhZ2VtZW50LCBhbmQgYXV0b21hdGluZyBtYW51YWwgd29ya2Zsb3dz4oCUbGlrZSBhcHBvaW50bWVudCBzY2hlZHVsaW5nIGFuZCBldmVudCBzaWduLXVwc+KAlHNvIHlvdXIgdGVhbSBzcGVuZHMgbGVzcyB0aW1lIG9uIGFkbWluIGFuZCBtb3JlIG9uIGFkdmlzaW5nLiBXZSBhbHNvIG9mZmVyIHRhaWxvcmVkIGNvbnN1bHRpbmcgdG8gc3RyZWFtbGluZSBldmVudCBtYW5hZ2VtZW50IGFuZCBvcHRpbWl6ZSBmZWVkYmFjayBjb2xsZWN0aW9uLgpXb3VsZCB5b3UgYmUgb3BlbiB0byBhIHF1aWNrIGNhbGwgbmV4dCB3ZWVrIHRvIGV4cGxvcmUgaG93IHdlIGNhbiBzdXBwb3J0IHlvdXIgbWlzc2lvbj8="
$ws.Range("J3").Value = ""
$ws.Range("E3").Value = 4.6
$ws.Range("F3").Value = 19
$ws.Range("M3").Value = 2.26

# ----- Row 4 -----
Set-CellB64 "A4" "Q2hJSkU5NGpvRDFSd29rUmFIX1gyVlhSRU5r"
Set-CellB64 "B4" "TkogU21hbGwgQnVzaW5lc3MgRGV2ZWxvcG1lbnQgQ2VudGVyIGF0IE5KQ1UgU2Nob29sIG9mIEJ1c2luZXNz"
Set-CellB64 "C4" "Y29ycG9yYXRlX29mZmljZSwgcG9pbnRfb2ZfaW50ZXJlc3QsIGVzdGFibGlzaG1lbnQ="
Set-CellB64 "D4" "KDIwMSkgMjAwLTIxNTY="
Set-CellB64 "H4" "aHR0cHM6Ly93d3cuc2JkY2h1ZHNvbmNvdW50eS5jb20v"
Set-CellB64 "I4" "aHR0cHM6Ly9tYXBzLmdvb2dsZS5jb20vP2NpZD0xNTY0MTIzMTY3MjUxNDkzNjY4MCZnX21wPUNpZG5iMjluYkdVdWJXRndjeTV3YkdGalpYTXVkakV1VUd4aFkyVnpMbE5sWVhKamFGUmxlSFFRQWhnRUlBQQ=="
Set-CellB64 "K4" "SGlnaGx5IHJlY29tbWVuZCEgSSB0b29rIGEgY291cGxlIHdvcmtzaG9wcyBlYXJsaWVyIHRoaXMgeWVhciB0aGF0IGhhdmUgaGVscGVkIHRyZW1lbmRvdXNseS4gSSB0b29rIHRoZSBjb3Vyc2Ugd2hpY2ggYXNzaXN0ZWQgaW4gZXN0YWJsaXNoaW5nIG15IGxsYyBhbmQgdGF4IGlkIGNyZWF0aW9uIGFuZCBhbHNvIHRoZSBidXNpbmVzcyBwbGFuIHdvcmsgc2hvcC4gUGxlbnR5IG9mIGJlbmVmaWNpYWwgd2ViaW5hcnMgYXJlIGJlaW5nIG9mZmVyZWQgZXZlbiBpbiB0aGUgcGFuZGVtaWMgd2hpY2ggaXMgZ3JlYXQuIEdyYXRlZnVsIGZvciB0aGVzZSBhZmZvcmRhYmxlIHJlc291cmNlcyBhcyBhIG5ldyBidXNpbmVzcyBvd25lciEgfCBUaGV5IG9mZmVyIHVuYmVsaWV2YWJsZSBhbmQgZXh0cmVtZWx5IHJlc291cmNlZnVsIGNvdXJzZXMgZm9yIGZyZWUgZm9yIGEgIHBlcnNvbiB0byBiZSBpbmZvcm0gYW5kIHdpbiBhdCB0aGVpciBidXNpbmVzcy4gIEFsbCBwZXJzb25uZWwgaXMgd2VsbCAgZWR1Y2F0ZWQsIGluZm9ybWVkIGFuZCB3YW50cyB0byBzZWUgeW91IHdpbi4gIEV4dHJlbWVseSBoZWxwZnVsLCBrbm93bGVkZ2VhYmxlIGFuZCB0aG91Z2h0ZnVsIGNvdW5zZWxvcnMgd2hvIGd1aWRlIHlvdSBpbiB0aGUgcmlnaHQgZGlyZWN0aW9uIGFzIGFuIGVudHJlcHJlbmV1ciAgd2l0aCBubyBrbm93bGVkZ2Ugb2YgYnVzaW5lc3MgKGJ1dCBvbmx5IHRoZSB0aG91Z2h0IHRvIHdpbikhICBUaGVyZSBpcyBubyBmcmVlIGhlbHAgZm9yIHRob3NlIHdobyB3YW50IHRvIGV4Y2VsIGF0IHN0YXJpbmcgdGhlaXIgb3duIHZlbnR1cmVzIGFuZCB3YW50IGV2ZXJ5b25lIGFyb3VuZCB0aGVtIHRvIHNoaW5lIGFzIHdlbGwsIG9ubHkgdGhpcyBwcm9ncmFtISB8IEkgaGF2ZSBhdHRlbmRlZCBtYW55IHRyYWluaW5ncyB3aXRoIE5KU0RDLiAgVGhlIGluZm9ybWF0aW9uIEkgcmVjZWl2ZWQgd2FzIHZpdGFsIHRvIG1lIGJ1aWxkaW5nIHRoZSBmb3VuZGF0aW9uIGZvciBTdGVmYW5pZSBCbHVlIENvbnN1bHRpbmcuIFRoZXkgd2VyZSBrbm93bGVkZ2VhYmxlIGluIGV2ZXJ5IGFzcGVjdCBvZiBidXNpbmVzcyBmcm9tIGZpbmFuY2UgdG8gbWFya2V0aW5nIGFuZCBzYWxlcy4gV2hldGhlciB5b3UgYXJlIHByb3ZpZGluZyBhIHNlcnZpY2Ugb3Igc2VsbGluZyBhIHByb2R1Y3QgdGhlIGluZm9ybWF0aW9uIGNvdmVycyBldmVyeSBxdWVzdGlvbiB5b3UgbWF5IGhhdmUgYWJvdXQgc3RhcnRpbmcgYSBuZXcgYnVzaW5lc3MuICBBTkQgdGhleSBtYWtlIHRoZW1zZWx2ZXMgYXZhaWxhYmxlIHRvIHlvdSBhZnRlciB0aGUgdHJhaW5pbmdzLg=="
Set-CellB64 "L4" "c2JkY0BuamN1LmVkdQ=="
Set-CellB64 "N4" "SGVyZSBhcmUgdGhlICoqVUkgaW1wcm92ZW1lbnRzKiogaWRlbnRpZmllZCBmb3IgdGhlIHdlYnNpdGUgc2NyZWVuc2hvdDoKLS0tCiMjIyAqKjEuIE91dGRhdGVkIFZpc3VhbCBTdHlsZSoqCi0gVGhlIGRlc2lnbiBkb2VzIG5vdCBmb2xsb3cgdGhlIGxhdGVzdCBtaW5pbWFsaXN0aWMgYW5kIOKAnGNsZWFu4oCdIFVJIHRyZW5kczsgaXQgYXBwZWFycyBzb21ld2hhdCBkYXRlZCB3aXRoIGhlYXZ5IGRyb3Agc2hhZG93cywgZ3JhZGllbnRzLCBhbmQgZGVuc2UgY29udGVudCBibG9ja3MuCi0gVGhlcmUgaXMgYSByZWxpYW5jZSBvbiBsYXJnZSwgYm9sZCB0ZXh0IHJhdGhlciB0aGFuIHZpc3VhbCBoaWVyYXJjaHkgYW5kIHdoaXRlIHNwYWNlIGZvciBlbXBoYXNpcy4KIyMjICoqMi4gSW5jb25zaXN0ZW50IEJyYW5kaW5nICYgQ29sb3IgUGFsZXR0ZSoqCi0gVGhlIHBhbGV0dGUgaXMgcGF0cmlvdGljIChyZWQsIHdoaXRlLCBibHVlKSwgYnV0IHRoZSBzaGFkZXMgYW5kIHNhdHVyYXRpb24gdmFyeSBiZXR3ZWVuIHNlY3Rpb25zLCBtYWtpbmcgdGhlIGJyYW5kaW5nIGZlZWwgaW5jb25zaXN0ZW50LgotIFRoZSDigJxMZWFybiBNb3Jl4oCdIGFuZCDigJxSZWdpc3RlciBmb3IgU3VjY2Vzc+KAnSBidXR0b25zIHVzZSBkaWZmZXJlbnQgdG9uZXMgb2YgcmVkLgotIFZpc3VhbCBoaWVyYXJjaHkgYmV0d2VlbiBzZWN0aW9ucyBpcyBtdWRkbGVkIGR1ZSB0byBpbmNvbnNpc3RlbnQgYmFja2dyb3VuZCBjb2xvcnMgYW5kIGNhbGwtdG8tYWN0aW9uIHN0eWxlcy4KIyMjICoqMy4gVHlwb2dyYXBoeSBJc3N1ZXMqKgotIFRvbyBtYW55IGZvbnQgd2VpZ2h0cywgc2l6ZXMsIGFuZCBjb2xvcnMgY3JlYXRlIGEgY2x1dHRlcmVkIGxvb2suIEZvciBleGFtcGxlLCB0aGUgbWl4IG9mIHZpYnJhbnQgd2hpdGUvZ3JleSBhbmQgYm9sZCBmb250cyBpbiBoZWFkaW5ncyBpcyB2aXN1YWxseSBqYXJyaW5nLgotIFBhcmFncmFwaHMgYXJlIGp1c3RpZmllZCBpbmNvbnNpc3RlbnRseSBhbmQgY291bGQgYmVuZWZpdCBmcm9tIG1vcmUgbGluZSBzcGFjaW5nIGZvciBpbXByb3ZlZCByZWFkYWJpbGl0eS4KIyMjICoqNC4gTGF5b3V0ICYgU3BhY2luZyoqCi0gVGhlcmUgaXMgbm90IGVub3VnaCBwYWRkaW5nL3NwYWNlIGJldHdlZW4gbWFueSBzZWN0aW9ucywgbGVhZGluZyB0byBhIGNyb3dkZWQgZmVlbGluZy4KLSBUaGUgYmxvY2tzIGZvciDigJxTdGFydOKAnSwg4oCcR3Jvd+KAnSwg4oCcRXhwYW5k4oCdIGFyZSB0aGUgc2FtZSB3aWR0aCBidXQgY291bGQgdXNlIG1vcmUgcGFkZGluZyBhbmQgbGVzcyB0ZXh0LCBwZXJoYXBzIHdpdGggaWNvbnMgdG8gaW1wcm92ZSBhdC1hLWdsYW5jZSB1bmRlcnN0YW5kaW5nLgotIFRoZSBzZWFyY2ggZmllbGQgdW5kZXIg4oCcVHJhaW5pbmcgJiBFdmVudHPigJ0gaXMgdmVyeSBzcXVlZXplZCBhbmQgdmlzdWFsbHkgdW5iYWxhbmNlZC4KIyMjICoqNS4gQWNjZXNzaWJpbGl0eS9Vc2FiaWxpdHkqKgotIEluc3VmZmljaWVudCBjb2xvciBjb250cmFzdCBpbiBzb21lIGFyZWFzIChsaWdodCBncmF5IGJhY2tncm91bmRzIHdpdGggd2hpdGUgb3IgbGlnaHQgdGV4dCkuCi0gU29tZSB0ZXh0IGJsb2NrcyBhcmUgaGFyZCB0byByZWFkIGR1ZSB0byBzaXplIG9yIGxhY2sgb2Ygc3BhY2luZy4KLSBCdXR0b25zIGNvdWxkIGJlIGxhcmdlciBhbmQgaGF2ZSBhIGNsZWFyZXIgZm9jdXMgc3RhdGUgZm9yIGtleWJvYXJkIG5hdmlnYXRpb24uCi0gVGhlIHBsYWNlaG9sZGVyIHRleHQgaW4gdGhlIG5ld3NsZXR0ZXIgZm9ybSBpcyBhbG1vc3QgdGhlIHNhbWUgY29sb3IgYXMgdGhlIGJhY2tncm91bmQsIHJlZHVjaW5nIHJlYWRhYmlsaXR5LgojIyMgKio2LiBWaXN1YWwgQ29uc2lzdGVuY3kqKgotIFRoZSBuYXZpZ2F0aW9uIGJhciBtZW51LCBwb3B1cCwgYW5kIHNlYXJjaCBpY29uIHN0eWxlcyBkaWZmZXIgZnJvbSB0aGUgcmVzdCBvZiB0aGUgc2l0ZS4KLSBUaGUg4oCcSGV5LCBXZWxjb21lIeKAnSBzZWN0aW9uIGZlZWxzIHNlcGFyYXRlZCBzdHlsaXN0aWNhbGx5IGFuZCBsYWNrcyB2aXN1YWwgaW50ZWdyYXRpb24gd2l0aCB0aGUgbWFpbiBicmFuZC4KLSBVbmV2ZW4gdXNlIG9mIGRyb3Agc2hhZG93cyAoZS5nLiwgYmVoaW5kIHRoZSBtb25pdG9yIG1vY2t1cCkgZmVlbHMgaW5jb25zaXN0ZW50IHdpdGggZmxhdCwgbW9kZXJuIFVJIHRyZW5kcy4KIyMjICoqNy4gSW1hZ2VyeSoqCi0gSGVybyBpbWFnZSBpcyBuaWNlIGJ1dCBsb3ctY29udHJhc3QgYW5kIGRvZXNu4oCZdCBjcmVhdGUgYSBzdHJvbmcgYnJhbmQgZmVlbGluZzsgY29uc2lkZXIgaGlnaGVyIGNvbnRyYXN0IG9yIGJsZW5kZWQgb3ZlcmxheXMgZm9yIGxlZ2liaWxpdHkuCi0gU3RvY2sgaW1hZ2VyeSAodGhlIGNvbXB1dGVyIG1vbml0b3IpIGRvZXNu4oCZdCBhZGQgdmFsdWUg4oCUIGN1c3RvbSBvciBtb3JlIHJlbGV2YW50IHZpc3VhbHMgd291bGQgZW5oYW5jZSBjcmVkaWJpbGl0eS4KIyMjICoqOC4gRm9vdGVyKioKLSBUaGUgZm9vdGVyIGFyZWEgaXMgZGVuc2Ugd2l0aCBzbWFsbCB0ZXh0IGxpbmtzOyB0aGVzZSBjb3VsZCBiZSBtYWRlIG1vcmUgcmVhZGFibGUgYW5kIHNlcGFyYXRlZCB3aXRoIGluY3JlYXNlZCBsaW5lIGhlaWdodCBvciBhcHByb3ByaWF0ZSBzcGFjaW5nLgotIFNvY2lhbCBpY29ucyBsYWNrIGVub3VnaCBwYWRkaW5nIGFuZCBjb250cmFzdCBmb3IgcHJvbWluZW5jZS4KLS0tCiMjICoqU3VtbWFyeSBvZiBSZWNvbW1lbmRhdGlvbnMqKgotIEFkb3B0IGEgbW9yZSBtaW5pbWFsLCBtb2Rlcm4gbG9vayB3aXRoIGxhcmdlciB3aGl0ZSBzcGFjZXMsIGZsYXR0ZXIgY29sb3JzLCBhbmQgZmV3ZXIgZm9udCBzdHlsZXMuCi0gVXNlIGEgY29oZXNpdmUgY29sb3IgcGFsZXR0ZSB0aHJvdWdob3V0OyBhbGlnbiBhbGwgQ1RBIGJ1dHRvbnMgdG8gYSBzaW5nbGUgc3R5bGUuCi0gSW5jcmVhc2Ugc3BhY2luZyBiZXR3ZWVuIHNlY3Rpb25zIGFuZCBhdm9pZCBjcm93ZGluZyB0ZXh0IGJsb2Nrcy4KLSBFbnN1cmUgYWxsIHRleHQgYW5kIFVJIGNvbXBvbmVudHMgaGF2ZSBzdWZmaWNpZW50IGNvbnRyYXN0IGZvciBhY2Nlc3NpYmlsaXR5LgotIFJlcGxhY2Ugb3IgcmVtb3ZlIG91dGRhdGVkIHN0b2NrIGltYWdlczsgdXNlIGNvbnNpc3RlbnQgaWNvbm9ncmFwaHkuCi0gUmVkZXNpZ24gdGhlIG5ld3NsZXR0ZXIgZm9ybSBmb3IgY2xhcml0eSBhbmQgbWFrZSBpbnB1dCBmaWVsZHMgbW9yZSBwcm9taW5lbnQuCi0gUmVmaW5lIHRoZSBuYXZpZ2F0aW9uIGFuZCBmb290ZXIgZm9yIGJldHRlciB1c2FiaWxpdHkgYW5kIHVzZXIgZXhwZXJpZW5jZS4KLS0tCioqT3ZlcmFsbCwgdGhlIHNpdGUgaXMgZnVuY3Rpb25hbCBidXQgZmVlbHMgb3V0ZGF0ZWQgd2l0aCBzZXZlcmFsIFVJL1VYIHByYWN0aWNlcyB0aGF0IGRvbuKAmXQgYWxpZ24gd2l0aCBtb2Rlcm4gd2ViIHN0YW5kYXJkcyBhbmQgYWNjZXNzaWJpbGl0eSB0cmVuZHMuKio="
Set-CellB64 "O4" "KipDbGllbnQgQnJpZWY6KioKVGhlIFNtYWxsIEJ1c2luZXNzIERldmVsb3BtZW50IENlbnRlciAoU0JEQykgYXQgTmV3IEplcnNleSBDaXR5IFVuaXZlcnNpdHkgaXMgYSBub24tcHJvZml0LCBnb3Zlcm5tZW50LWZ1bmRlZCByZXNvdXJjZSBjZW50ZXIgc3VwcG9ydGluZyBlbnRyZXByZW5ldXJzIGFuZCBzbWFsbCBidXNpbmVzcyBvd25lcnMgaW4gSHVkc29uIENvdW50eSwgTmV3IEplcnNleS4gVGhleSBvZmZlciBuby1jb3N0IG9uZS1vbi1vbmUgY291bnNlbGluZywgdHJhaW5pbmcsIGFuZCBhZHZpc29yeSBzZXJ2aWNlcyB0byBoZWxwIGJ1c2luZXNzZXMgc3RhcnQsIGdyb3csIGFuZCBleHBhbmQsIHdpdGggZXhwZXJ0aXNlIHNwYW5uaW5nIGFjY2VzcyB0byBjYXBpdGFsLCBzdHJhdGVnaWMgcGxhbm5pbmcsIG1hcmtldGluZywgSFIsIGFuZCBtb3JlLiBBcyBhIHBpdm90YWwgaHViIGZvciBsb2NhbCBidXNpbmVzcyBlZHVjYXRpb24gYW5kIHN1cHBvcnQsIFNCREMgY291bGQgYmVuZWZpdCBmcm9tIGF1dG9tYXRpbmcgYXBwb2ludG1lbnQgYm9va2luZywgY2VudHJhbGl6aW5nIGV2ZW50IGFuZCByZXNvdXJjZSBtYXRlcmlhbHMsIGFuZCByZWRlc2lnbmluZyBpdHMgd2Vic2l0ZSBmb3IgZWFzaWVyIG5hdmlnYXRpb24gYW5kIG1vYmlsZSByZXNwb25zaXZlbmVzcy4gU3RyZWFtbGluaW5nIGludGFrZSBhbmQgZm9sbG93LXVwIHByb2Nlc3NlcyB3b3VsZCBhbHNvIG1heGltaXplIHRoZWlyIGFkdmlzb3Jz4oCZIGltcGFjdCwgd2hpbGUgaW1wcm92ZW1lbnRzIHRvIHRoZSBuZXdzbGV0dGVyIHNpZ251cCBleHBlcmllbmNlIGNvdWxkIGJvb3N0IGNvbW11bml0eSBlbmdhZ2VtZW50IGFuZCBpbmZvcm1hdGlvbiBkaXNzZW1pbmF0aW9uLg=="
Set-CellB64 "P4" "KipCdXNpbmVzcyBQYWluIFBvaW50IFJlcG9ydCoqCi0tLQojIyMgUGFpbiBQb2ludCAxOiBPdXRkYXRlZCwgSW5jb25zaXN0ZW50LCBhbmQgVW5pbnR1aXRpdmUgV2Vic2l0ZSBEZXNpZ24gIAotICoqSW1wYWN0OioqICAKICBUaGUgY3VycmVudCB1c2VyIGV4cGVyaWVuY2Ugb24gdGhlIHdlYnNpdGUgaXMgYmVsb3cgbW9kZXJuIHN0YW5kYXJkcywgd2hpY2ggY2FuIGxlYWQgdG8gdmlzaXRvciBmcnVzdHJhdGlvbiwgcmVkdWNlZCBlbmdhZ2VtZW50IHdpdGggcmVzb3VyY2VzLCBhbmQgZGlmZmljdWx0eSBhY2Nlc3NpbmcgZXNzZW50aWFsIHNlcnZpY2VzLiBQb3RlbnRpYWwgY2xpZW50cyBtaWdodCBwZXJjZWl2ZSB0aGUgY2VudGVyIGFzIGxlc3MgcHJvZ3Jlc3NpdmUgb3IgdHJ1c3R3b3J0aHksIHBvdGVudGlhbGx5IGh1cnRpbmcgc2lnbi11cHMgZm9yIGNvdW5zZWxpbmcsIHdvcmtzaG9wcywgYW5kIG5ld3NsZXR0ZXJzLiBBY2Nlc3NpYmlsaXR5IGlzc3VlcyBjb3VsZCBhbHNvIGxpbWl0IG91dHJlYWNoIHRvIHVzZXJzIHdpdGggZGlzYWJpbGl0aWVzLCB3aGlsZSBtb2JpbGUgdW5mcmllbmRsaW5lc3MgaGFtcGVycyBlbmdhZ2VtZW50IGZyb20gYnVzeSBlbnRyZXByZW5ldXJzIG9uLXRoZS1nby4KLSAqKkhvdyBXZSBDYW4gSGVscDoqKiAgCiAgKldlYnNpdGUgUmVkZXNpZ246KiBXZSBvZmZlciBhIGdyb3VuZC11cCBtb2Rlcm5pemF0aW9uIG9mIHRoZSB3ZWJzaXRl4oCZcyBkZXNpZ24sIGFwcGx5aW5nIGNsZWFuLCBhY2Nlc3NpYmxlLCBhbmQgbW9iaWxlLWZpcnN0IFVJL1VYIHN0YW5kYXJkcy4gVGhpcyBpbmNsdWRlcyBhIHVuaWZpZWQgY29sb3IgcGFsZXR0ZSwgc3Ryb25nIHZpc3VhbCBoaWVyYXJjaHksIG9wdGltYWwgc3BhY2luZywgYWNjZXNzaWJsZSBlbGVtZW50cywgYW5kIGNvaGVzaXZlIGJyYW5kaW5nIHRvIHJlaW5mb3JjZSB0cnVzdCBhbmQgY3JlZGliaWxpdHkg4oCUIGVuc3VyaW5nIGEgZnJpY3Rpb25sZXNzIGpvdXJuZXkgZnJvbSBsYW5kaW5nIHBhZ2UgdG8gdXNlciBhY3Rpb24uCi0tLQojIyMgUGFpbiBQb2ludCAyOiBNYW51YWwvRnJhZ21lbnRlZCBPcGVyYXRpb25zIGZvciBBcHBvaW50bWVudCBCb29raW5nIGFuZCBJbnRha2UgIAotICoqSW1wYWN0OioqICAKICBNYW51YWwgb3IgdW5pbnRlZ3JhdGVkIHByb2Nlc3NlcyBmb3IgYXBwb2ludG1lbnQgc2NoZWR1bGluZywgY2xpZW50IG9uYm9hcmRpbmcsIGFuZCBwb3N0LXNlc3Npb24gZm9sbG93LXVwIGFyZSBwcm9uZSB0byBlcnJvcnMgYW5kIGJvdHRsZW5lY2tzLCBhYnNvcmJpbmcgdmFsdWFibGUgc3RhZmYgYW5kIGFkdmlzb3IgdGltZS4gVGhpcyBkaW1pbmlzaGVzIHRoZSBjZW50ZXLigJlzIGNhcGFjaXR5IHRvIHNlcnZlIG1vcmUgY2xpZW50cywgY2F1c2VzIHBvdGVudGlhbCBsb3N0IG9wcG9ydHVuaXRpZXMsIGFuZCBzbG93cyB0aGUgb3ZlcmFsbCBjbGllbnQgam91cm5leS4KLSAqKkhvdyBXZSBDYW4gSGVscDoqKiAgCiAgKkF1dG9tYXRpb24gb2YgT3BlcmF0aW9uczoqIFdlIGNhbiBpbXBsZW1lbnQgc2VjdXJlLCB1c2VyLWZyaWVuZGx5IGF1dG9tYXRpb24gc29sdXRpb25zIGZvciBhcHBvaW50bWVudCBib29raW5nLCBpbnRha2UgZm9ybXMsIGFuZCBmb2xsb3ctdXAgY29tbXVuaWNhdGlvbnMuIEJ5IGNlbnRyYWxpemluZyBhbmQgZGlnaXRpemluZyB0aGVzZSB3b3JrZmxvd3MsIHRoZSBTQkRDIHdpbGwgc2F2ZSBzdGFmZiBob3VycywgaW5jcmVhc2Ugc2NoZWR1bGluZyBhY2N1cmFjeSwgYW5kIHByb3ZpZGUgYSBzZWFtbGVzcyBleHBlcmllbmNlIGZvciBjbGllbnRzIGZyb20gZmlyc3QgcG9pbnQgb2YgY29udGFjdCB0aHJvdWdoIG9uZ29pbmcgc3VwcG9ydC4KLS0tCiMjIyBQYWluIFBvaW50IDM6IERpc2pvaW50ZWQgQWNjZXNzIHRvIFJlc291cmNlIE1hdGVyaWFscywgVHJhaW5pbmcgRXZlbnRzLCBhbmQgTmV3c2xldHRlcnMgIAotICoqSW1wYWN0OioqICAKICBWYWx1YWJsZSBpbmZvcm1hdGlvbiBpcyBoYXJkZXIgZm9yIGVudHJlcHJlbmV1cnMgdG8gZmluZCBvciBhY2Nlc3MsIHJlc3VsdGluZyBpbiBsb3dlciBldmVudCBhdHRlbmRhbmNlLCBtaXNzZWQgb3Bwb3J0dW5pdGllcyBmb3IgZWR1Y2F0aW9uLCBhbmQgdW5kZXJ1dGlsaXphdGlvbiBvZiBhdmFpbGFibGUgcmVzb3VyY2VzLiBBIGNvbmZ1c2luZyBvciBjdW1iZXJzb21lIG5ld3NsZXR0ZXIgc2lnbnVwIHByb2Nlc3MgcmVkdWNlcyBlbmdhZ2VtZW50IGFuZCBkaW1pbmlzaGVzIHRoZSBjZW50ZXLigJlzIGFiaWxpdHkgdG8gaW5mb3JtIGFuZCBtb2JpbGl6ZSBpdHMgY29tbXVuaXR5LgotICoqSG93IFdlIENhbiBIZWxwOioqICAKICAqUHJvY2VzcyBDb25zdWx0aW5nICYgV2Vic2l0ZSBJbnRlZ3JhdGlvbjoqIFdlIHdpbGwgc3RyZWFtbGluZSB0aGUgd2Vic2l0ZeKAmXMgcmVzb3VyY2UgY2VudGVyLCBldmVudCBtYW5hZ2VtZW50LCBhbmQgbmV3c2xldHRlciBmdW5jdGlvbnMuIFRoaXMgaW5jbHVkZXMgZGV2ZWxvcGluZyBhIGNlbnRyYWxpemVkIHJlc291cmNlIGh1YiwgaW1wcm92aW5nIGZpbHRlcmluZy9zZWFyY2ggZm9yIHRyYWluaW5ncywgYW5kIHJldmFtcGluZyB0aGUgbmV3c2xldHRlciBzaWdudXAgZm9yIG1heGltdW0gdmlzaWJpbGl0eSBhbmQgZWFzZS4gVGhlc2UgY2hhbmdlcyBoZWxwIGJvb3N0IHBhcnRpY2lwYXRpb24gcmF0ZXMsIG1ha2Uga25vd2xlZGdlIG1vcmUgYWNjZXNzaWJsZSwgYW5kIGRlZXBlbiBvbmdvaW5nIHJlbGF0aW9uc2hpcHMgd2l0aCB0aGUgc21hbGwgYnVzaW5lc3MgY29tbXVuaXR5LgotLS0KIyMjIFBhaW4gUG9pbnQgNDogSW5jb25zaXN0ZW50IEJyYW5kaW5nIGFuZCBXZWFrIFZpc3VhbCBJZGVudGl0eSBPbmxpbmUgIAotICoqSW1wYWN0OioqICAKICBWYXJpYWJpbGl0eSBpbiBjb2xvciBzY2hlbWVzLCBidXR0b24gc3R5bGVzLCBhbmQgaWNvbm9ncmFwaHkgY2FuIGRpbHV0ZSBicmFuZCByZWNvZ25pdGlvbiBhbmQgcmVkdWNlIHRoZSBwZXJjZWl2ZWQgcHJvZmVzc2lvbmFsaXNtIG9mIHRoZSBvcmdhbml6YXRpb24uIFRoaXMgd2Vha2VucyBlbW90aW9uYWwgY29ubmVjdGlvbiBhbmQgbWFrZXMgaXQgaGFyZGVyIGZvciBTQkRDIHRvIHN0YW5kIG91dCBhcyBhIGxlYWRlciBhbmQgdHJ1c3RlZCBhZHZpc2VyIGluIHRoZSBleWVzIG9mIGxvY2FsIGVudHJlcHJlbmV1cnMgYW5kIHBhcnRuZXJzLgotICoqSG93IFdlIENhbiBIZWxwOioqICAKICAqQnJhbmQgQWxpZ25tZW50IFRocm91Z2ggUmVkZXNpZ246KiBPdXIgdGVhbSBjYW4gYXVkaXQgYW5kIHVuaWZ5IHRoZSBkaWdpdGFsIHZpc3VhbCBpZGVudGl0eSwgY3JlYXRpbmcgYSBzdHlsZSBndWlkZSB0aGF0IGVuc3VyZXMgY29uc2lzdGVudCBicmFuZGluZywgdG9uZSwgYW5kIHVzYWJpbGl0eSBhY3Jvc3Mgd2ViIHBhZ2VzIGFuZCBjb21tdW5pY2F0aW9uIGNoYW5uZWxzLiBBIG1vZGVybiwgcmVsaWFibGUgbG9vayBpbmNyZWFzZXMgdHJ1c3QsIHN0cmVuZ3RoZW5zIHlvdXIgYnJhbmQgaW4gdGhlIG1pbmRzIG9mIHZpc2l0b3JzLCBhbmQgc3VwcG9ydHMgb3JnYW5pemF0aW9uYWwgZ3Jvd3RoLgotLS0KIyMjIFBhaW4gUG9pbnQgNTogTGFjayBvZiBTeXN0ZW1hdGljIE1lYXN1cmVtZW50IGFuZCBDb250aW51b3VzIEltcHJvdmVtZW50ICAKLSAqKkltcGFjdDoqKiAgCiAgV2l0aG91dCBpbnRlZ3JhdGVkIGRpZ2l0YWwgdG91Y2hwb2ludHMgYW5kIHJlcG9ydGluZywgZ2F1Z2luZyB0aGUgZWZmZWN0aXZlbmVzcyBvZiBvdXRyZWFjaCBpbml0aWF0aXZlcyBhbmQgY2xpZW50IGVuZ2FnZW1lbnQgZWZmb3J0cyBpcyBkaWZmaWN1bHQuIE1pc3Npbmcgb3IgaW5jb21wbGV0ZSBkYXRhIGNhbiByZXN1bHQgaW4gbWlzc2VkIG9wcG9ydHVuaXRpZXMgdG8gaW1wcm92ZSBvciBleHBhbmQgaGlnaC1pbXBhY3Qgc2VydmljZXMsIGFuZCBoaW5kZXJzIHJlcG9ydGluZyB0byBmdW5kZXJzIG9yIHN0YWtlaG9sZGVycy4KLSAqKkhvdyBXZSBDYW4gSGVscDoqKiAgCiAgKkF1dG9tYXRpb24gJiBBbmFseXRpY3MgSW50ZWdyYXRpb246KiBXZSBjYW4gZW5hYmxlIHdlYnNpdGUgYW5kIHByb2Nlc3MgYXV0b21hdGlvbiB3aXRoIGJ1aWx0LWluIGFuYWx5dGljcyDigJQgdHJhY2tpbmcgc2lnbi11cHMsIHJlc291cmNlIGRvd25sb2FkcywgZXZlbnQgcGFydGljaXBhdGlvbiBhbmQgbW9yZS4gVGhpcyBkYXRhIHdpbGwgZW1wb3dlciBsZWFkZXJzaGlwIHRvIG1ha2UgaW5mb3JtZWQgZGVjaXNpb25zLCBvcHRpbWl6ZSBwcm9ncmFtcywgYW5kIGRlbW9uc3RyYXRlIHZhbHVlIHRvIGZ1bmRlcnMgYW5kIHBhcnRuZXJzLgotLS0KKipTdW1tYXJ5IGZvciBPdXRyZWFjaDoqKiAgCldoaWxlIFNCREMgYXQgTkpDVSBpcyBoaWdobHkgdmFsdWVkIGJ5IGl0cyBjb21tdW5pdHksIHRoZXJlIGlzIHN1YnN0YW50aWFsIG9wcG9ydHVuaXR5IHRvIG11bHRpcGx5IGl0cyBpbXBhY3QgdGhyb3VnaCBtb2Rlcm5pemF0aW9uIG9mIGl0cyBvbmxpbmUgcHJlc2VuY2UsIGF1dG9tYXRpb24gb2Ygc3RhZmYtaW50ZW5zaXZlIHByb2Nlc3NlcywgYW5kIGludGVncmF0aW9uIG9mIGNlbnRyYWwgcmVzb3VyY2VzLiBPdXIgcHJvdmVuIGV4cGVydGlzZSBjYW4gaGVscCBTQkRDIHN0cmVhbWxpbmUgb3BlcmF0aW9ucywgaW1wcm92ZSBjbGllbnQgYW5kIHN0YWZmIGV4cGVyaWVuY2UsIGFuZCBzdHJlbmd0aGVuIHRoZSBvcmdhbml6YXRpb27igJlzIHN0YW5kaW5nIGFzIGEgZm9yd2FyZC1sb29raW5nIHBpbGxhciBmb3IgTmV3IEplcnNleeKAmXMgc21hbGwgYnVzaW5lc3Nlcy4="
Set-CellB64 "Q4" "U3ViamVjdDogVW5sb2NrIEdyZWF0ZXIgSW1wYWN0IGZvciBTQkRDIGF0IE5KQ1UKSGkgW05hbWVdLApJ4oCZbSBpbXByZXNzZWQgYnkgdGhlIFNCREPigJlzIGNvbW1pdG1lbnQgdG8gSHVkc29uIENvdW50eeKAmXMgZW50cmVwcmVuZXVycywgYnV0IEkgbm90aWNlZCB0aGUgd2Vic2l0ZSBhbmQgaW50YWtlIHByb2Nlc3MgY291bGQgYmUgaG9sZGluZyB5b3VyIHRlYW0gYmFjay4gQSBtb2Rlcm4sIG1vYmlsZS1mcmllbmRseSByZWRlc2lnbiBhbmQgYXV0b21hdGVkIGJvb2tpbmcgc3lzdGVtIHdvdWxkIGhlbHAgeW91ciBhZHZpc29ycyBzZXJ2ZSBtb3JlIHNtYWxsIGJ1c2luZXNzZXPigJR3aGlsZSBhIHN0cmVhbWxpbmVkIHJlc291cmNlIGFuZCBldmVudCBodWIgYm9vc3RzIGNvbW11bml0eSBlbmdhZ2VtZW50LiBBcmUgeW91IG9wZW4gdG8gYSBxdWljayBjYWxsIHRvIGRpc2N1c3MgaG93IHdlIGNhbiBoZWxwIFNCREMgZGVsaXZlciBldmVuIGdyZWF0ZXIgaW1wYWN0Pw=="
$ws.Range("E4").Value = 4.8
$ws.Range("F4").Value = 81
$ws.Range("M4").Value = 2.19

# ----- Row 5 -----
Set-CellB64 "A5" "Q2hJSmFRNzY4MzNpd29rUlZwdWE3X2NNLU00"
Set-CellB64 "B5" "UFJPRkVTU0lPTkFMIEJVU0lORVNTIFNUUkFURUdJRVMsIExMQw=="
Set-CellB64 "C5" "Y29uc3VsdGFudCwgcG9pbnRfb2ZfaW50ZXJlc3QsIGVzdGFibGlzaG1lbnQ="
Set-CellB64 "D5" "KDIwMSkgNDI1LTgzMDA="
Set-CellB64 "H5" "aHR0cDovL3d3dy5wYnN0cmF0ZWdpZXMuY29tLz91dG1fc291cmNlPWdvb2dsZSZ1dG1fbWVkaXVtPXdpeF9nb29nbGVfYnVzaW5lc3NfcHJvZmlsZSZ1dG1fY2FtcGFpZ249MTAxNTE3OTExNjc4MzYyMTEwNDU="
Set-CellB64 "I5" "aHR0cHM6Ly9tYXBzLmdvb2dsZS5jb20vP2NpZD0xNDkxMzY4NDQyNTA1MzczOTg2MiZnX21wPUNpZG5iMjluYkdVdWJXRndjeTV3YkdGalpYTXVkakV1VUd4aFkyVnpMbE5sWVhKamFGUmxlSFFRQWhnRUlBQQ=="
Set-CellB64 "K5" "SSdtIHRocmlsbGVkIHRvIHNoYXJlIG15IDUtc3RhciBleHBlcmllbmNlIHdpdGggRGVubmlzIEJvbmFndXJhISBPdmVyIHRoZSBwYXN0IDUgeWVhcnMsIERlbm5pcyBoYXMgYmVlbiBhbiBpbnN0cnVtZW50YWwgcGFydG5lciBpbiBteSBidXNpbmVzcyBqb3VybmV5LCBndWlkaW5nIG1lIGZyb20gJDEwMCwwMDAgaW4gbXkgZmlyc3QgeWVhciB0byAkMS43IG1pbGxpb24gaW4gYW5udWFsIHJldmVudWUgKHllYXIgNSkuIFdoYXQgc2V0cyBEZW5uaXMgYXBhcnQgaXMgbm90IGp1c3QgaGlzIGV4Y2VwdGlvbmFsIGJ1c2luZXNzIGV4cGVydGlzZSwgYnV0IGhpcyBnZW51aW5lIGNhcmUgYW5kIG1lbnRvcnNoaXAuIEhlJ3MgbW9yZSB0aGFuIGp1c3QgbXkgY29hY2ggLSBoZSdzIGEgdHJ1c3RlZCBmcmllbmQgYW5kIGFkdmlzb3Igd2hvJ3MgYWx3YXlzIHRoZXJlIHRvIG9mZmVyIG1lIHZhbHVhYmxlIGluc2lnaHRzLCBjZWxlYnJhdGUgc3VjY2Vzc2VzLCBhbmQgcHJvdmlkZSBndWlkYW5jZSBkdXJpbmcgY2hhbGxlbmdpbmcgdGltZXMuCkRlbm5pcyBoYXMgYSB1bmlxdWUgYWJpbGl0eSB0byBwaW5wb2ludCBhcmVhcyBmb3IgaW1wcm92ZW1lbnQgYW5kIGhpZ2hsaWdodCBzdHJlbmd0aHMsIHByb3ZpZGluZyBhY3Rpb25hYmxlIGFkdmljZSB0aGF0IGRyaXZlcyByZXN1bHRzLiBIaXMgc3RyYXRlZ2llcyBhcmUgcHJhY3RpY2FsLCBhbmQgaGlzIHN1cHBvcnQgaXMgdW53YXZlcmluZy4gSSBoaWdobHkgcmVjb21tZW5kIERlbm5pcyB0byBhbnlvbmUgc2Vla2luZyBhIGJlbmVmaWNpYWwgYnVzaW5lc3MgY29hY2ggd2hvIGNhbiBvZmZlciBhIGZhaXIgYW5kIGFmZm9yZGFibGUgZ3VpZGFuY2UsIGFzIHdlbGwgYXMgYSB0cnVzdGVkIHBhcnRuZXIgdG8gbmF2aWdhdGUgYm90aCB0cml1bXBocyBhbmQgc2V0YmFja3Mgd2l0aC4KNSBzdGFycyBpc24ndCBlbm91Z2ggLSBEZW5uaXMgaGFzIHRydWx5IHRyYW5zZm9ybWVkIG15IGJ1c2luZXNzIGFuZCBteSB3aG9sZSBsaWZlIQpEYXZpZCBQZWVycwpQcmVzaWRlbnQKQ3V0dGluZyBFZGdlIENvbnN1bHRpbmcgR3JvdXAgSW5jLiB8IEdyZWF0IGV4cGVyaWVuY2Ugd29ya2luZyB3aXRoIFBCUy4gIERlbm5pcyBoYXMgaGFkIGEgaHVnZSBpbXBhY3Qgb24gdGhlIHN1Y2Nlc3Mgb2YgbXkgcGx1bWJpbmcgJiBoZWF0aW5nIGNvbXBhbnkuCkhlIGhhcyBjZXJ0YWlubHkgYnJvdWdodCBDTEFSSVRZIHRvIFdIQVQgTUFUVEVSUyBNT1NUISEgfCBQcm9mZXNzaW9uYWwgQnVzaW5lc3MgU3RyYXRlZ2llcyBpcyBhIHRydWx5IGEgZm9yd2FyZCB0aGlua2luZyBmaXJtIGR1ZSB0byB0aGVpciBrbm93bGVkZ2UgYW5kIGV4cGVyaWVuY2UgaW4gZGVhbGluZyB3aXRoIHRvZGF5J3MgY29tcGxpY2F0ZWQgd29ya3BsYWNlLgpJJ3ZlIGtub3duIERlbm5pcyBCb25hZ3VyYSBmb3IgeWVhcnMgYnV0IGl0IHdhc24ndCB1bnRpbCBJIHdhcyB1bmVtcGxveWVkIHRoYXQgSSB0cnVseSBsZWFybmVkIG9mIHRoZSBncmVhdCBidXNpbmVzcyBzZW5zZSBoZSBmb3N0ZXJlZCBpbiBwZXJzb25hbCBkZXZlbG9wbWVudC4gIEhlIGxvb2tlZCBhdCBteSByZXN1bWUgd2hpY2ggd2FzIHByb2Zlc3Npb25hbGx5IGRvbmUgYW5kIGluIGhpcyBjaGFybWluZyB3YXkgc2FpZCBpdCB3YXNuJ3Qgd2hhdCB0b2RheSdzIGhpcmluZyBtYW5hZ2VyIHdhcyBsb29raW5nIGZvciBpbiB0aGUgNDUgc2Vjb25kcyB0aGV5J2Qgc3BlbmQgb24gdGhlIGZpcnN0IGdsYW5jZSEgIEhlIGFkdmlzZWQgdGhlIGNoYW5nZXMgdGhhdCBJIHNob3VsZCBtYWtlIGFuZCB0byB0cmFuc2ZlciB0aGF0IHRvIG15IExpbmtlZEluIHByb2ZpbGUgYXMgd2VsbC4gIEkgc3RhcnRlZCByZWNlaXZpbmcgbW9yZSBoaXRzIGltbWVkaWF0ZWx5IGFuZCB3YXMgYmFjayB0byB3b3JrIHdpdGhpbiBhIG1vbnRoLgpJIHRha2UgZ3JlYXQgcGxlYXN1cmUgdGhhbmtpbmcgRGVubmlzIGV2ZXJ5IHRpbWUgSSBzZWUgaGltIQpKb2huIE1jQ2F1bGV5ClNlbmlvciBQcm9qZWN0IE1hbmFnZXIKR2xvYmFsIEJhbmtpbmcgTGlxdWlkaXR5IFJlZ3VsYXRpb25z"
Set-CellB64 "L5" "ZGJvbmFndXJhQHBic3RyYXRlZ2llcy5jb20="
Set-CellB64 "N5" "KipVSSBJbXByb3ZlbWVudHM6KioKMS4gKipPdXRkYXRlZCBWaXN1YWwgU3R5bGU6KioKICAgLSBUaGUgd2Vic2l0ZeKAmXMgZGVzaWduIGZlZWxzIGRhdGVkLCBsYWNraW5nIHRoZSBtb2Rlcm4gbWluaW1hbGlzbSBhbmQgcG9saXNoIHNlZW4gaW4gY29udGVtcG9yYXJ5IHdlYnNpdGVzIChlLmcuLCBtb3JlIHdoaXRlc3BhY2UsIGNhcmQtYmFzZWQgbGF5b3V0cywgc3VidGxlIGdyYWRpZW50cykuCiAgIC0gQmx1cnJlZCBzdG9jayBpbWFnZXMgYXJlIG92ZXJ1c2VkIGFuZCBkZXRyYWN0IGZyb20gdmlzdWFsIGFwcGVhbCBhbmQgY2xhcml0eS4gTW9kZXJuIHdlYnNpdGVzIGZhdm9yIGhpZ2gtcXVhbGl0eSBvciBicmFuZC1zcGVjaWZpYyBpbWFnZXJ5LgoyLiAqKkJyYW5kIENvaGVzaW9uOioqCiAgIC0gVGhlIGxvZ28gYW5kIGNvbG9yIHNjaGVtZSAoYmxhY2ssIGJsdWUsIGFuZCB3aGl0ZSkgYXJlIGNvbnNpc3RlbnQsIGJ1dCB0aGUgaW1hZ2VyeSBkb2VzIG5vdCByZWluZm9yY2UgdGhlIGJyYW5kIGlkZW50aXR5LiBNb3JlIGN1c3RvbWl6ZWQgb3Igb24tYnJhbmQgaW1hZ2VzIHdvdWxkIGhlbHAuCiAgIC0gRm9udCB1c2FnZSBpcyBub3QgZGlzdGluY3RpdmUgb3IgbW9kZXJuOyBpbmNvcnBvcmF0aW5nIGEgdW5pcXVlLCBsZWdpYmxlIHR5cGVmYWNlIGNhbiBoZWxwIHdpdGggYnJhbmQgcmVjb2duaXRpb24uCjMuICoqTGF5b3V0ICYgU3BhY2luZzoqKgogICAtIENvbnRlbnQgc2VjdGlvbnMgYXJlIHZlcnkgZGVuc2UsIGVzcGVjaWFsbHkgaW4gbW9iaWxlIHZpZXcuIEFkZCBtb3JlIHBhZGRpbmcvc3BhY2luZyBiZXR3ZWVuIHNlY3Rpb25zIGFuZCBjb250ZW50IGJsb2NrcyBmb3IgcmVhZGFiaWxpdHkgYW5kIHZpc3VhbCBlYXNlLgogICAtIFRoZSBjb250YWN0IGZvcm0gYWJydXB0bHkgY2hhbmdlcyB0byBhIHN0YXJrIHdoaXRlIGJhY2tncm91bmQsIGJyZWFraW5nIHRoZSB2aXN1YWwgZmxvdy4gTWFpbnRhaW4gY29uc2lzdGVudCBiYWNrZ3JvdW5kIGNvbG9yIG9yIHN1YnRsZSBzZXBhcmF0b3IuCjQuICoqVHlwb2dyYXBoeToqKgogICAtIEZvbnQgc2l6ZXMgZm9yIGJvZHkgYW5kIGhlYWRlciB0ZXh0IGxhY2sgaGllcmFyY2h5LiBVc2UgbGFyZ2VyIGhlYWRpbmdzLCBtb3JlIHByb25vdW5jZWQgc3ViaGVhZGluZ3MsIGFuZCBiYWxhbmNlZCBib2R5IHRleHQgZm9yIGVhc2llciBzY2FubmluZy4KICAgLSBUZXh0IG92ZXIgaW1hZ2VzIGxhY2tzIHN1ZmZpY2llbnQgY29udHJhc3QsIGVzcGVjaWFsbHkgaW4gdGhlIGhlcm8gc2VjdGlvbiAod2hpdGUgdGV4dCBvbiBhIGJyaWdodCBpbWFnZSkuIEFwcGx5IG92ZXJsYXlzIG9yIGRhcmtlciBpbWFnZSBmaWx0ZXJzIHRvIGluY3JlYXNlIHJlYWRhYmlsaXR5Lgo1LiAqKkltYWdlcyAmIFZpc3VhbHM6KioKICAgLSBSZXBlYXRlZCBibHVycnkgaW1hZ2VzIGRpbHV0ZSBwcm9mZXNzaW9uYWxpc207IHVzZSBzaGFycCwgaGlnaC1yZXMgaW1hZ2VzIG9yIHN1YnRsZSBpbGx1c3RyYXRpb25zL2ljb25zIHRoYXQgbWF0Y2ggdGhlIGJyYW5k4oCZcyBtZXNzYWdlLgogICAtIENvbnNpZGVyIG1vZGVybiBkZXNpZ24gZWxlbWVudHMgc3VjaCBhcyBzb2Z0IGRyb3Agc2hhZG93cywgcm91bmRlZCBjb3JuZXJzLCBvciBtaWNyby1hbmltYXRpb25zIGZvciBidXR0b25zIGFuZCBpbWFnZXMuCjYuICoqQnV0dG9ucyAmIENhbGxzLXRvLUFjdGlvbjoqKgogICAtICJSZWFkIE1vcmUiIGJ1dHRvbnMgYXJlIGdlbmVyaWM7IHVzZSBtb3JlIGRlc2NyaXB0aXZlIGxhYmVscyB0byBndWlkZSB1c2VycyAoZS5nLiwgIkV4cGxvcmUgU2VydmljZXMsIiAiVmlldyBBc3Nlc3NtZW50IERldGFpbHMiKS4KICAgLSBCdXR0b25zIGxhY2sgdmlzdWFsIGVtcGhhc2lzIChzaXplL2NvbG9yL3NoYXBlKTsgbWFrZSB0aGVtIHN0YW5kIG91dCBtb3JlIHdpdGggY29udHJhc3QgYW5kIHN1YnRsZSBhbmltYXRpb24uCjcuICoqQWNjZXNzaWJpbGl0eToqKgogICAtIFRleHQgY29udHJhc3QgKGVzcGVjaWFsbHkgdGV4dCBvdmVyIGltYWdlcykgaXMgaW5zdWZmaWNpZW50IGZvciBhY2Nlc3NpYmlsaXR5LiBVc2UgV0NBRy1jb21wbGlhbnQgY29udHJhc3QgcmF0aW9zLgogICAtIE5hdmlnYXRpb24vdGV4dCBzaXplIGFuZCBidXR0b24gZWxlbWVudHMgc2hvdWxkIGJlIGxhcmdlci9tb3JlIHRvdWNoLWZyaWVuZGx5IG9uIG1vYmlsZS4KICAgLSBDb25zaWRlciBza2lwIG5hdmlnYXRpb24gYW5kIEFSSUEgcm9sZXMgZm9yIGltcHJvdmVkIGFjY2Vzc2liaWxpdHkuCjguICoqRm9vdGVyOioqCiAgIC0gVGhlIGZvb3RlciBmZWVscyBkaXNjb25uZWN0ZWQ7IGluY3JlYXNlIHBhZGRpbmcsIG9yZ2FuaXplIGxpbmtzIGludG8gY29sdW1ucywgYW5kIHVzZSBpY29ucyBmb3Igc29jaWFsIG1lZGlhLgogICAtIFRoZSB1c2Ugb2YgdGhyZWUgdmVyeSBkaWZmZXJlbnQgYmFja2dyb3VuZCBjb2xvcnMgKGJsYWNrLCB3aGl0ZSwgYmx1ZSkgaGFtcGVycyBjb2hlc2lvbi4gVW5pZnkgYmFja2dyb3VuZHMgZm9yIGEgc21vb3RoZXIgZmxvdy4KOS4gKipVc2FiaWxpdHk6KioKICAgLSBBZGQgc3RpY2t5L2ZpeGVkIG5hdmlnYXRpb24gZm9yIGVhc2llciBhY2Nlc3MgdG8gc2VjdGlvbnMuCiAgIC0gVGhlIHRlc3RpbW9uaWFscyBzZWN0aW9uIHNob3VsZCB1c2UgY29uc2lzdGVudCBjYXJkIGhlaWdodHMgYW5kIGNvbnNpZGVyIGEgY2Fyb3VzZWwvc2xpZGVyIGZvciBtb2Rlcm4gcHJlc2VudGF0aW9uLgoxMC4gKipNb2Rlcm4gVHJlbmRzOioqCiAgICAtIEFkZCBzdWJ0bGUgYW5pbWF0aW9ucyBvciB0cmFuc2l0aW9ucyB0byBpbWFnZXMvYnV0dG9ucyBmb3IgYSBtb2Rlcm4gdG91Y2guCiAgICAtIEludGVncmF0ZSBzdWJ0bGUgZ3JhZGllbnRzIG9yIGdsYXNzbW9ycGhpc20gYXMgYSBiYWNrZ3JvdW5kIG9yIGhpZ2hsaWdodCB0byBjYXRjaCB1cCB3aXRoIGN1cnJlbnQgVUkgdHJlbmRzLgotLS0KKipTdW1tYXJ5OioqICAKV2hpbGUgdGhlIHdlYnNpdGUgaGFzIGEgY2xlYXIgc3RydWN0dXJlIGFuZCBiYXNpYyBjb2hlc2lvbiwgaXQgcmVxdWlyZXMgc2lnbmlmaWNhbnQgbW9kZXJuaXphdGlvbiBpbiBpbWFnZXJ5LCB2aXN1YWwgZmxvdywgYWNjZXNzaWJpbGl0eSwgYW5kIGludGVyYWN0aXZpdHkgdG8gY29tcGV0ZSB3aXRoIGN1cnJlbnQgVUkvVVggc3RhbmRhcmRzIGFuZCBwcm9kdWNlIGEgbW9yZSBwb2xpc2hlZCwgcHJvZmVzc2lvbmFsIGV4cGVyaWVuY2Uu"
Set-CellB64 "O5" "KipDbGllbnQgQnJpZWY6KioKUHJvZmVzc2lvbmFsIEJ1c2luZXNzIFN0cmF0ZWdpZXMsIExMQywgbGVkIGJ5IERlbm5pcyBCb25hZ3VyYSBpbiBCZXJnZW4gQ291bnR5LCBOSiwgaXMgYSBidXNpbmVzcyBjb2FjaGluZyBhbmQgY29uc3VsdGluZyBmaXJtIHNwZWNpYWxpemluZyBpbiBzdXBwb3J0aW5nIGJ1c2luZXNzIG93bmVycywgZXhlY3V0aXZlcywgYW5kIHNhbGVzIHByb2Zlc3Npb25hbHMuIERlbm5pcyBpcyBhIENlcnRpZmllZCBQcm9mZXNzaW9uYWwgQnVzaW5lc3MgQ29hY2ggd2l0aCBvdmVyIDMwIHllYXJzJyBsZWFkZXJzaGlwIGV4cGVyaWVuY2UsIG9mZmVyaW5nIHRhaWxvcmVkIGNvYWNoaW5nLCBleGVjdXRpdmUgbWVudG9yaW5nLCBzYWxlcyB0ZWFtIGRldmVsb3BtZW50LCBhbmQgYmVoYXZpb3JhbCBhc3Nlc3NtZW50cyAoRXh0ZW5kZWQgRElTQ8KuKSB0byBoZWxwIGNsaWVudHMgZ2FpbiBjbGFyaXR5LCBpbXByb3ZlIHBlcmZvcm1hbmNlLCBhbmQgcmVzb2x2ZSBjb21tb24gb3JnYW5pemF0aW9uYWwgY2hhbGxlbmdlcyBsaWtlIHN0YWduYW50IHNhbGVzLCBwcm9maXRhYmlsaXR5LCBhbmQgc3RhZmYgdHVybm92ZXIuIFF1aWNrIHdpbnMgY291bGQgaW5jbHVkZSBhdXRvbWF0aW5nIGFwcG9pbnRtZW50IHNjaGVkdWxpbmcgYW5kIGFzc2Vzc21lbnQgZGVsaXZlcnksIHJldmFtcGluZyB0aGUgd2Vic2l0ZSBmb3IgY2xlYW5lciBuYXZpZ2F0aW9uIGFuZCBjbGVhcmVyIGNhbGxzLXRvLWFjdGlvbiAoZS5nLiwgIkJvb2sgYSBGcmVlIFN0cmF0ZWd5IFNlc3Npb24iKSwgYW5kIHN5c3RlbWF0aXppbmcgcmVjdXJyaW5nIGNsaWVudCBmZWVkYmFjayByZXF1ZXN0cy4gR2l2ZW4gdGhlIHN0cm9uZyBmb2N1cyBvbiBjbGFyaXR5IGFuZCBhY3Rpb25hYmxlIHN0cmF0ZWdpZXMsIGFuIG9wZXJhdGlvbnMvcHJvY2VzcyBhdWRpdCBjb3VsZCBmdXJ0aGVyIGlkZW50aWZ5IGJvdHRsZW5lY2tzIGluIGNsaWVudCBvbmJvYXJkaW5nIG9yIGZvbGxvdy11cCwgd2hpbGUgdGFyZ2V0ZWQgd2ViIHVwZGF0ZXMgd291bGQgcmVpbmZvcmNlIHRoZWlyIGV4cGVydCBwb3NpdGlvbmluZyBhbmQgZW5oYW5jZSBsZWFkIGNhcHR1cmUu"
Set-CellB64 "P5" "KipCdXNpbmVzcyBQYWluIFBvaW50IFJlcG9ydCoqICAKKkZvciBQcm9mZXNzaW9uYWwgQnVzaW5lc3MgU3RyYXRlZ2llcywgTExDIChQQlMpIOKAkyBEZW5uaXMgQm9uYWd1cmEqCi0tLQoqKlBhaW4gUG9pbnQgMTogT3V0ZGF0ZWQgV2Vic2l0ZSBEZXNpZ24gYW5kIFVzZXIgRXhwZXJpZW5jZSoqICAKLSAqKkltcGFjdDoqKiAgCiAgVGhlIHdlYnNpdGUgZmVlbHMgb3V0ZGF0ZWQgYW5kIGxhY2tzIG1vZGVybiBwb2xpc2gsIHdpdGggaXNzdWVzIGluIGltYWdlcnksIHR5cG9ncmFwaHksIHNwYWNpbmcsIGNvbnRyYXN0LCBjYWxscy10by1hY3Rpb24sIGFuZCBhY2Nlc3NpYmlsaXR5LiBUaGlzIHJpc2tzIHVuZGVybWluaW5nIHRoZSBmaXJt4oCZcyBjcmVkaWJpbGl0eSwgZGV0ZXJyaW5nIHBvdGVudGlhbCBjbGllbnRzLCBhbmQgcmVkdWNpbmcgY29udmVyc2lvbiByYXRlcyBmb3Iga2V5IGdvYWxzIChlLmcuLCBib29raW5nIHN0cmF0ZWd5IHNlc3Npb25zIG9yIGFzc2Vzc21lbnRzKS4gVGhlIGxhY2sgb2YgbW9kZXJuIGRlc2lnbiBhbmQgYWNjZXNzaWJpbGl0eSBjYW4gbWFrZSBQQlMgYXBwZWFyIGxlc3MgcmVsZXZhbnQgb3IgcHJvZmVzc2lvbmFsIGNvbXBhcmVkIHRvIGNvbXBldGl0b3JzLCBpbXBhY3RpbmcgbGVhZCBnZW5lcmF0aW9uIGFuZCBmaXJzdCBpbXByZXNzaW9ucy4KLSAqKkhvdyBXZSBDYW4gSGVscDoqKiAgCiAgKipXZWJzaXRlIFJlZGVzaWduOioqIE91ciB3ZWIgZGV2ZWxvcG1lbnQgdGVhbSB3aWxsIGltcGxlbWVudCBhIG1vZGVybiwgYnJhbmQtY29oZXNpdmUsIGFuZCBhY2Nlc3NpYmxlIHJlZGVzaWdu4oCUdXRpbGl6aW5nIGhpZ2gtcXVhbGl0eSBpbWFnZXMsIGltcHJvdmVkIHR5cG9ncmFwaHksIGNvbnNpc3RlbnQgYnJhbmRpbmcsIGNsZWFyIG5hdmlnYXRpb24sIGFuZCBlbmdhZ2luZyBjYWxscy10by1hY3Rpb24gKGUuZy4sIOKAnEJvb2sgYSBGcmVlIFN0cmF0ZWd5IFNlc3Npb27igJ0pLiBVcGRhdGVzIHdpbGwgbWFrZSB0aGUgd2Vic2l0ZSBtb2JpbGUtZnJpZW5kbHksIHZpc3VhbGx5IGFwcGVhbGluZywgYW5kIG9wdGltaXplZCBmb3IgbGVhZCBjYXB0dXJlLCByZWluZm9yY2luZyBQQlPigJlzIGV4cGVydCBwb3NpdGlvbmluZyBhbmQgc3VwcG9ydGluZyBidXNpbmVzcyBncm93dGguCi0tLQoqKlBhaW4gUG9pbnQgMjogSW5lZmZpY2llbnQgYW5kIE1hbnVhbCBCdXNpbmVzcyBQcm9jZXNzZXMqKiAgCi0gKipJbXBhY3Q6KiogIAogIFJlY3VycmluZyBhZG1pbmlzdHJhdGl2ZSB0YXNrcyBzdWNoIGFzIGFwcG9pbnRtZW50IHNjaGVkdWxpbmcsIGFzc2Vzc21lbnQgZGVsaXZlcnksIGFuZCBjbGllbnQgZmVlZGJhY2sgcmVxdWVzdHMgYXJlIGxpa2VseSBoYW5kbGVkIG1hbnVhbGx5LCBjcmVhdGluZyBvcGVyYXRpb25hbCBib3R0bGVuZWNrcywgd2FzdGluZyBzdGFmZiB0aW1lLCBpbmNyZWFzaW5nIGVycm9yIHJpc2ssIGFuZCBzbG93aW5nIGRvd24gY2xpZW50IG9uYm9hcmRpbmcgb3IgZm9sbG93LXVwLiBUaGlzIGluZWZmaWNpZW5jeSBjYW4gdW5kZXJtaW5lIHRoZSBwZXJzb25hbGl6ZWQgY29hY2hpbmcgZXhwZXJpZW5jZSB0aGF0IERlbm5pcyBvZmZlcnMsIHJlc3RyaWN0IGNhcGFjaXR5IGZvciBncm93dGgsIGFuZCBmcnVzdHJhdGUgY2xpZW50cyBleHBlY3RpbmcgYSBzZWFtbGVzcywgbW9kZXJuIGV4cGVyaWVuY2UuCi0gKipIb3cgV2UgQ2FuIEhlbHA6KiogIAogICoqQnVzaW5lc3MgUHJvY2VzcyBBdXRvbWF0aW9uOioqIFdlIGNhbiBhdXRvbWF0ZSBrZXkgd29ya2Zsb3dz4oCUaW5jbHVkaW5nIGNsaWVudCBzY2hlZHVsaW5nLCBhc3Nlc3NtZW50L3JlcG9ydCBkZWxpdmVyeSwgY2xpZW50IGZlZWRiYWNrIGNvbGxlY3Rpb24sIGFuZCByb3V0aW5lIGNvbW11bmljYXRpb27igJR1c2luZyBwbGF0Zm9ybXMgbGlrZSBDYWxlbmRseSwgWmFwaWVyLCBvciBjdXN0b20gYXV0b21hdGlvbi4gVGhpcyBzYXZlcyB0aW1lLCByZWR1Y2VzIGh1bWFuIGVycm9yLCBjcmVhdGVzIGEgc21vb3RoZXIgY2xpZW50IGpvdXJuZXksIGFuZCBhbGxvd3MgRGVubmlzIHRvIGZvY3VzIG9uIHZhbHVlLWFkZGluZyB3b3JrIHJhdGhlciB0aGFuIHJlcGV0aXRpdmUgdGFza3MuCi0tLQoqKlBhaW4gUG9pbnQgMzogTGltaXRlZCBTeXN0ZW1hdGljIEFwcHJvYWNoIHRvIENvbGxlY3RpbmcgYW5kIExldmVyYWdpbmcgQ2xpZW50IEZlZWRiYWNrKiogIAotICoqSW1wYWN0OioqICAKICBXaGlsZSByZXZpZXdzIGluZGljYXRlIGhpZ2ggc2F0aXNmYWN0aW9uLCBmZWVkYmFjayBjb2xsZWN0aW9uIHNlZW1zIGFkIGhvYy4gTGFjayBvZiBhIHN5c3RlbWF0aWMgcHJvY2VzcyBjYW4gbWVhbiBtaXNzZWQgaW5zaWdodHMgZm9yIGJ1c2luZXNzIGRldmVsb3BtZW50LCBsZXNzIHNvY2lhbCBwcm9vZiBvbiBkaWdpdGFsIGNoYW5uZWxzLCBhbmQgYW4gdW5kZXJ1dGlsaXphdGlvbiBvZiBjbGllbnQgYWR2b2NhY3kgZm9yIG1hcmtldGluZy4gVGhpcyByZXN1bHRzIGluIGZld2VyIGNvbXBlbGxpbmcgdGVzdGltb25pYWxzIGFuZCBhIHdlYWtlciBkZW1vbnN0cmF0aW9uIG9mIHRoZSBmaXJt4oCZcyBpbXBhY3QgdG8gcHJvc3BlY3RpdmUgY2xpZW50cy4KLSAqKkhvdyBXZSBDYW4gSGVscDoqKiAgCiAgKipPcmdhbml6YXRpb25hbC9Qcm9jZXNzIENvbnN1bHRpbmc6KiogV2XigJlsbCBkZXZlbG9wIGEgc3RydWN0dXJlZCBmZWVkYmFjayBwcm9jZXNz4oCUYXV0b21hdGluZyByZWd1bGFyIGZvbGxvdy11cCByZXF1ZXN0cyBmb3IgdGVzdGltb25pYWxzLCBOZXQgUHJvbW90ZXIgU2NvcmUgKE5QUyksIGFuZCBzdXJ2ZXlpbmcgY2xpZW50cyBhdCBrZXkgbWlsZXN0b25lcy4gVGhpcyBub3Qgb25seSBpbmNyZWFzZXMgdGhlIHZvbHVtZSBvZiB2YWx1YWJsZSBjbGllbnQgZmVlZGJhY2sgdG8gZGlzcGxheSBvbiB0aGUgcmV2YW1wZWQgd2Vic2l0ZSBhbmQgbWFya2V0aW5nIG1hdGVyaWFscywgYnV0IGFsc28gcHJvdmlkZXMgRGVubmlzIHdpdGggYWN0aW9uYWJsZSBpbnNpZ2h0cyB0byByZWZpbmUgc2VydmljZSBkZWxpdmVyeS4KLS0tCioqUGFpbiBQb2ludCA0OiBVbmRlci1vcHRpbWl6ZWQgTGVhZCBHZW5lcmF0aW9uIGFuZCBDb252ZXJzaW9uIEZ1bm5lbCoqICAKLSAqKkltcGFjdDoqKiAgCiAgR2VuZXJpYyBidXR0b25zLCB1bmNsZWFyIGNhbGxzLXRvLWFjdGlvbiwgbGFjayBvZiBsZWFkIG1hZ25ldHMgKGxpa2UgZnJlZSBzdHJhdGVneSBzZXNzaW9ucyksIGFuZCBhIG5vbi1kaXN0aW5jdGl2ZSBicmFuZCBwcmVzZW5jZSBvbiB0aGUgd2Vic2l0ZSBtYXkgbGVhZCB0byBtaXNzZWQgb3Bwb3J0dW5pdGllcyB0byBjYXB0dXJlIGludGVyZXN0IGFuZCBxdWFsaWZ5IGxlYWRzLiBUaGlzIGRpcmVjdGx5IGltcGFjdHMgdGhlIHNhbGVzIHBpcGVsaW5lLCBsb3dlcmluZyB0aGUgbnVtYmVyIG9mIHF1YWxpdHkgbGVhZHMgYW5kIG1ha2luZyBtYXJrZXRpbmcgc3BlbmQgbGVzcyBlZmZlY3RpdmUuCi0gKipIb3cgV2UgQ2FuIEhlbHA6KiogIAogICoqU2FsZXMgT3BlcmF0aW9ucyBBdXRvbWF0aW9uICYgU3RyYXRlZ2ljIFdlYnNpdGUgVXBkYXRlczoqKiBXZeKAmWxsIG9wdGltaXplIGNhbGxzLXRvLWFjdGlvbiAoQ1RBcyksIGNyZWF0ZSBjb21wZWxsaW5nIGxlYWQgbWFnbmV0cyAoZS5nLiwgYXNzZXNzbWVudCBwcmV2aWV3cywgc3RyYXRlZ3kgc2Vzc2lvbiBib29raW5ncyksIGFuZCBzeXN0ZW1hdGl6ZSBsZWFkIGNhcHR1cmUgd2l0aCBpbnRlZ3JhdGVkIGZvcm1zIGFuZCBDUk0gc3luY2luZy4gVGhpcyBlbnN1cmVzIHRoYXQgd2Vic2l0ZSB2aXNpdG9ycyBhcmUgbW9yZSBsaWtlbHkgdG8gZW5nYWdlLCBib29rIHNlc3Npb25zLCBhbmQgZW50ZXIgdGhlIHNhbGVzIHBpcGVsaW5l4oCUaW1wcm92aW5nIGJvdGggY29udmVyc2lvbiByYXRlcyBhbmQgUk9JIGZyb20gZGlnaXRhbCB0cmFmZmljLgotLS0KKipQYWluIFBvaW50IDU6IFVudGFwcGVkIE9wcG9ydHVuaXRpZXMgaW4gRmluYW5jaWFsIGFuZCBPcGVyYXRpb25hbCBEYXRhIFJlcG9ydGluZyoqICAKLSAqKkltcGFjdDoqKiAgCiAgV2hpbGUgdGhlIHJldmlld3MgaGlnaGxpZ2h0IGV4Y2VwdGlvbmFsIGNvYWNoaW5nIHJlc3VsdHMsIHRoZXJl4oCZcyBubyBtZW50aW9uIG9mIHN0cnVjdHVyZWQgdHJhY2tpbmcvcmVwb3J0aW5nIG9uIGNsaWVudCBvciBpbnRlcm5hbCBLUElzLCBmaW5hbmNpYWwgbWV0cmljcywgb3Igb3BlcmF0aW9uYWwgcGVyZm9ybWFuY2UuIFdpdGhvdXQgdGhpcywgdGhlIGNvbXBhbnkgbWF5IGJlIG1pc3NpbmcgaW5zaWdodHMgZm9yIGRyaXZpbmcgcHJvZml0YWJpbGl0eSwgb3B0aW1pemluZyBzZXJ2aWNlIGRlbGl2ZXJ5LCBvciBjb21tdW5pY2F0aW5nIHZhbHVlIHRvIGNsaWVudHMgaW4gdGFuZ2libGUgd2F5cy4KLSAqKkhvdyBXZSBDYW4gSGVscDoqKiAgCiAgKipCb29ra2VlcGluZyAmIFBlcmZvcm1hbmNlIERhc2hib2FyZCBTZXR1cDoqKiBPdXIgdGVhbSBjYW4gaW1wbGVtZW50IHN0cmVhbWxpbmVkIGJvb2trZWVwaW5nIGFuZCBmaW5hbmNpYWwgcmVwb3J0aW5nLCBhcyB3ZWxsIGFzIGRhc2hib2FyZHMgZm9yIHRyYWNraW5nIGtleSBidXNpbmVzcyBhbmQgY2xpZW50IG1ldHJpY3MgKGdyb3d0aCwgZW5nYWdlbWVudCwgc2F0aXNmYWN0aW9uKS4gVGhpcyBlbmFibGVzIG1vcmUgaW5mb3JtZWQgc3RyYXRlZ2ljIGRlY2lzaW9ucyBmb3IgUEJTIGFuZCBlbmhhbmNlcyB0aGUgZmlybeKAmXMgYWJpbGl0eSB0byBkZW1vbnN0cmF0ZSBjbGllbnQgUk9JIGluIGZ1dHVyZSBlbmdhZ2VtZW50cy4KLS0tCioqU3VtbWFyeToqKiAgClByb2Zlc3Npb25hbCBCdXNpbmVzcyBTdHJhdGVnaWVzLCBMTEMgaXMgd2VsbCByZWdhcmRlZCBmb3IgcGVyc29uYWxpemVkLCByZXN1bHQtZHJpdmVuIGNvYWNoaW5nLCBidXQgZmFjZXMgY2xlYXIgcGFpbiBwb2ludHMgaW4gd2ViIHByZXNlbmNlLCBwcm9jZXNzIGVmZmljaWVuY3ksIHN5c3RlbWF0aWMgY2xpZW50IGZlZWRiYWNrLCBsZWFkIGdlbmVyYXRpb24sIGFuZCBkYXRhLWRyaXZlbiBvcGVyYXRpb25zLiBXZSByZWNvbW1lbmQgc3RhcnRpbmcgd2l0aCBhIGNvbXByZWhlbnNpdmUgd2Vic2l0ZSByZWRlc2lnbiBhbmQga2V5IHByb2Nlc3MgYXV0b21hdGlvbnMsIGZvbGxvd2VkIGJ5IG9wdGltaXphdGlvbiBvZiBmZWVkYmFjayBhbmQgbGVhZCBlbmdhZ2VtZW50IHByb2Nlc3NlcywgYW5kIHN1cHBvcnQgaW4gZmluYW5jaWFsIHJlcG9ydGluZyBmb3IgbWF4aW11bSBidXNpbmVzcyBpbXBhY3QuCi0tLQ=="
Set-CellB64 "Q5" "U3ViamVjdDogUXVpY2sgV2lucyBmb3IgUEJTOiBXZWJzaXRlLCBBdXRvbWF0aW9uICYgR3Jvd3RoCkhpIERlbm5pcywgIApJ4oCZbSBpbXByZXNzZWQgYnkgeW91ciBpbXBhY3Qgd2l0aCBQQlMgaGVscGluZyBsZWFkZXJzIGluIEJlcmdlbiBDb3VudHksIGJ1dCBub3RpY2VkIHlvdXIgd2Vic2l0ZSBhbmQgd29ya2Zsb3cgY291bGQgYmUgaG9sZGluZyB5b3UgYmFjay4gQSBtb2Rlcm4gc2l0ZSByZWRlc2lnbiB3aXRoIGNsZWFyZXIgY2FsbHMtdG8tYWN0aW9uIGNhbiBib29zdCBjcmVkaWJpbGl0eSBhbmQgbGVhZCBjYXB0dXJlLCB3aGlsZSBhdXRvbWF0aW5nIHRhc2tzIGxpa2Ugc2NoZWR1bGluZyBhbmQgYXNzZXNzbWVudCBkZWxpdmVyeSB3aWxsIGZyZWUgdXAgeW91ciB0aW1lIGZvciBjbGllbnRz4oCUbm90IGFkbWluLiBJIGNhbiBhbHNvIGhlbHAgeW91IHN5c3RlbWF0aXplIGZlZWRiYWNrIGNvbGxlY3Rpb24gdG8gc2hvd2Nhc2UgeW91ciBzdWNjZXNzIGFuZCBkcml2ZSBncm93dGguICAKV291bGQgeW91IGJlIG9wZW4gdG8gYSBxdWljayBjYWxsIG5leHQgd2VlayB0byBkaXNjdXNzIGFjdGlvbmFibGUgaW1wcm92ZW1lbnRzPyAgCkJlc3QsICAKW1lvdXIgTmFtZV0="
$ws.Range("E5").Value = 4.6
$ws.Range("F5").Value = 9
$ws.Range("M5").Value = 2.2

Write-Host "Edit complete. UsedRange:" $ws.UsedRange.Address()
